$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (person_id_nbr) updates
$ws.Range("E2").Value = 83258
$ws.Range("G2").Value = 41631.0801484542
$ws.Range("H2").Value = 41629.5
$ws.Range("I2").Value = 24036.59878645714
$ws.Range("J2").Value = 577758081.2211126
$ws.Range("K2").Value = 0
$ws.Range("P2").Value = 2
$ws.Range("R2").Value = 405180
$ws.Range("S2").Value = 83258
$ws.Range("T2").Value = '[''0'', ''1'', ''10'', ''100'', ''1000'', ''10000'', ''10001'', ''10002'', ''10003'', ''10004'', ''10005'', ''10006'', ''10007'', ''10008'', ''10009'', ''1001'', ''10010'', ''10011'', ''10012'', ''10013'', ''10014'', ''10015'', ''10016'', ''10017'', ''10018'', ''10019'', ''1002'', ''10020'', ''10021'', ''10022'', ''10023'', ''10024'', ''10025'', ''10026'', ''10027'', ''10028'', ''10029'', ''1003'', ''10030'', ''10031'', ''10032'', ''10033'', ''10034'', ''10035'', ''10036'', ''10037'', ''10038'', ''10039'', ''1004'', ''10040'', ''10041'', ''10042'', ''10043'', ''10044'', ''10045'', ''10046'', ''10047'', ''10048'', ''10049'', ''1005'', ''10050'', ''10051'', ''10052'', ''10053'', ''10054'', ''10055'', ''10056'', ''10057'', ''10058'', ''10059'', ''1006'', ''10060'', ''10061'', ''10062'', ''10063'', ''10064'', ''10065'', ''10066'', ''10067'', ''10068'', ''10069'', ''1007'', ''10070'', ''10071'', ''10072'', ''10073'', ''10074'', ''10075'', ''10076'', ''10077'', ''10078'', ''10079'', ''1008'', ''10080'', ''10081'', ''10082'', ''10083'', ''10084'', ''10085'', ''10086'', ''10087'', ''10088'', ''10089'', ''1009'', ''10090'', ''10091'', ''10092'', ''10093'', ''10094'', ''10095'', ''10096'', ''10097'', ''10098'', ''10099'', ''101'', ''1010'', ''10100'', ''10101'', ''10102'', ''10103'', ''10104'', ''10105'', ''10106'', ''10107'', ''10108'', ''10109'', ''1011'', ''10110'', ''10111'', ''10112'', ''10113'', ''10114'', ''10115'', ''10116'', ''10117'', ''10118'', ''10119'', ''1012'', ''10120'', ''10121'', ''10122'', ''10123'', ''10124'', ''10125'', ''10126'', ''10127'', ''10128'', ''10129'', ''1013'', ''10130'', ''10131'', ''10132'', ''10133'', ''10134'', ''10135'', ''10136'', ''10137'', ''10138'', ''10139'', ''1014'', ''10140'', ''10141'', ''10142'', ''10143'', ''10144'', ''10145'', ''10146'', ''10147'', ''10148'', ''10149'', ''1015'', ''10150'', ''10151'', ''10152'', ''10153'', ''10154'', ''10155'', ''10156'', ''10157'', ''10158'', ''10159'', ''1016'', ''10160'', ''10161'', ''10162'', ''10163'', ''10164'', ''10165'', ''10166'', ''10167'', ''10168'', ''10169'', ''1017'', ''10170'', ''10171'', ''10172'', ''10173'', ''10174'', ''10175'', ''10176'', ''10177'', ''10178'', ''10179'', ''1018'', ''10180'', ''10181'', ''10182'', ''10183'', ''10184'', ''10185'', ''10186'', ''10187'', ''10188'', ''10189'', ''1019'', ''10190'', ''10191'', ''10192'', ''10193'', ''10194'', ''10195'', ''10196'', ''10197'', ''10198'', ''10199'', ''102'', ''1020'', ''10200'', ''10201'', ''10202'', ''10203'', ''10204'', ''10205'', ''10206'', ''10207'', ''10208'', ''10209'', ''1021'', ''10210'', ''10211'', ''10212'', ''10213'', ''10214'', ''10215'', ''10216'', ''10217'', ''10218'', ''10219'', ''1022'', ''10220'', ''10221'', ''10222'', ''10223'', ''10224'', ''10225'', ''10226'', ''10227'', ''10228'', ''10229'', ''1023'', ''10230'', ''10231'', ''10232'', ''10233'', ''10234'', ''10235'', ''10236'', ''10237'', ''10238'', ''10239'', ''1024'', ''10240'', ''10241'', ''10242'', ''10243'', ''10244'', ''10245'', ''10246'', ''10247'', ''10248'', ''10249'', ''1025'', ''10250'', ''10251'', ''10252'', ''10253'', ''10254'', ''10255'', ''10256'', ''10257'', ''10258'', ''10259'', ''1026'', ''10260'', ''10261'', ''10262'', ''10263'', ''10264'', ''10265'', ''10266'', ''10267'', ''10268'', ''10269'', ''1027'', ''10270'', ''10271'', ''10272'', ''10273'', ''10274'', ''10275'', ''10276'', ''10277'', ''10278'', ''10279'', ''1028'', ''10280'', ''10281'', ''10282'', ''10283'', ''10284'', ''10285'', ''10286'', ''10287'', ''10288'', ''10289'', ''1029'', ''10290'', ''10291'', ''10292'', ''10293'', ''10294'', ''10295'', ''10296'', ''10297'', ''10298'', ''10299'', ''103'', ''1030'', ''10300'', ''10301'', ''10302'', ''10303'', ''10304'', ''10305'', ''10306'', ''10307'', ''10308'', ''10309'', ''1031'', ''10310'', ''10311'', ''10312'', ''10313'', ''10314'', ''10315'', ''10316'', ''10317'', ''10318'', ''10319'', ''1032'', ''10320'', ''10321'', ''10322'', ''10323'', ''10324'', ''10325'', ''10326'', ''10327'', ''10328'', ''10329'', ''1033'', ''10330'', ''10331'', ''10332'', ''10333'', ''10334'', ''10335'', ''10336'', ''10337'', ''10338'', ''10339'', ''1034'', ''10340'', ''10341'', ''10342'', ''10343'', ''10344'', ''10345'', ''10346'', ''10347'', ''10348'', ''10349'', ''1035'', ''10350'', ''10351'', ''10352'', ''10353'', ''10354'', ''10355'', ''10356'', ''10357'', ''10358'', ''10359'', ''1036'', ''10360'', ''10361'', ''10362'', ''10363'', ''10364'', ''10365'', ''10366'', ''10367'', ''10368'', ''10369'', ''1037'', ''10370'', ''10371'', ''10372'', ''10373'', ''10374'', ''10375'', ''10376'', ''10377'', ''10378'', ''10379'', ''1038'', ''10380'', ''10381'', ''10382'', ''10383'', ''10384'', ''10385'', ''10386'', ''10387'', ''10388'', ''10389'', ''1039'', ''10390'', ''10391'', ''10392'', ''10393'', ''10394'', ''10395'', ''10396'', ''10397'', ''10398'', ''10399'', ''104'', ''1040'', ''10400'', ''10401'', ''10402'', ''10403'', ''10404'', ''10405'', ''10406'', ''10407'', ''10408'', ''10409'', ''1041'', ''10410'', ''10411'', ''10412'', ''10413'', ''10414'', ''10415'', ''10416'', ''10417'', ''10418'', ''10419'', ''1042'', ''10420'', ''10421'', ''10422'', ''10423'', ''10424'', ''10425'', ''10426'', ''10427'', ''10428'', ''10429'', ''1043'', ''10430'', ''10431'', ''10432'', ''10433'', ''10434'', ''10435'', ''10436'', ''10437'', ''10438'', ''10439'', ''1044'', ''10440'', ''10441'', ''10442'', ''10443'', ''10444'', ''10445'', ''10446'', ''10447'', ''10448'', ''10449'', ''1045'', ''10450'', ''10451'', ''10452'', ''10453'', ''10454'', ''10455'', ''10456'', ''10457'', ''10458'', ''10459'', ''1046'', ''10460'', ''10461'', ''10462'', ''10463'', ''10464'', ''10465'', ''10466'', ''10467'', ''10468'', ''10469'', ''1047'', ''10470'', ''10471'', ''10472'', ''10473'', ''10474'', ''10475'', ''10476'', ''10477'', ''10478'', ''10479'', ''1048'', ''10480'', ''10481'', ''10482'', ''10483'', ''10484'', ''10485'', ''10486'', ''10487'', ''10488'', ''10489'', ''1049'', ''10490'', ''10491'', ''10492'', ''10493'', ''10494'', ''10495'', ''10496'', ''10497'', ''10498'', ''10499'', ''105'', ''1050'', ''10500'', ''10501'', ''10502'', ''10503'', ''10504'', ''10505'', ''10506'', ''10507'', ''10508'', ''10509'', ''1051'', ''10510'', ''10511'', ''10512'', ''10513'', ''10514'', ''10515'', ''10516'', ''10517'', ''10518'', ''10519'', ''1052'', ''10520'', ''10521'', ''10522'', ''10523'', ''10524'', ''10525'', ''10526'', ''10527'', ''10528'', ''10529'', ''1053'', ''10530'', ''10531'', ''10532'', ''10533'', ''10534'', ''10535'', ''10536'', ''10537'', ''10538'', ''10539'', ''1054'', ''10540'', ''10541'', ''10542'', ''10543'', ''10544'', ''10545'', ''10546'', ''10547'', ''10548'', ''10549'', ''1055'', ''10550'', ''10551'', ''10552'', ''10553'', ''10554'', ''10555'', ''10556'', ''10557'', ''10558'', ''10559'', ''1056'', ''10560'', ''10561'', ''10562'', ''10563'', ''10564'', ''10565'', ''10566'', ''10567'', ''10568'', ''10569'', ''1057'', ''10570'', ''10571'', ''10572'', ''10573'', ''10574'', ''10575'', ''10576'', ''10577'', ''10578'', ''10579'', ''1058'', ''10580'', ''10581'', ''10582'', ''10583'', ''10584'', ''10585'', ''10586'', ''10587'', ''10588'', ''10589'', ''1059'', ''10590'', ''10591'', ''10592'', ''10593'', ''10594'', ''10595'', ''10596'', ''10597'', ''10598'', ''10599'', ''106'', ''1060'', ''10600'', ''10601'', ''10602'', ''10603'', ''10604'', ''10605'', ''10606'', ''10607'', ''10608'', ''10609'', ''1061'', ''10610'', ''10611'', ''10612'', ''10613'', ''10614'', ''10615'', ''10616'', ''10617'', ''10618'', ''10619'', ''1062'', ''10620'', ''10621'', ''10622'', ''10623'', ''10624'', ''10625'', ''10626'', ''10627'', ''10628'', ''10629'', ''1063'', ''10630'', ''10631'', ''10632'', ''10633'', ''10634'', ''10635'', ''10636'', ''10637'', ''10638'', ''10639'', ''1064'', ''10640'', ''10641'', ''10642'', ''10643'', ''10644'', ''10645'', ''10646'', ''10647'', ''10648'', ''10649'', ''1065'', ''10650'', ''10651'', ''10652'', ''10653'', ''10654'', ''10655'', ''10656'', ''10657'', ''10658'', ''10659'', ''1066'', ''10660'', ''10661'', ''10662'', ''10663'', ''10664'', ''10665'', ''10666'', ''10667'', ''10668'', ''10669'', ''1067'', ''10670'', ''10671'', ''10672'', ''10673'', ''10674'', ''10675'', ''10676'', ''10677'', ''10678'', ''10679'', ''1068'', ''10680'', ''10681'', ''10682'', ''10683'', ''10684'', ''10685'', ''10686'', ''10687'', ''10688'', ''10689'', ''1069'', ''10690'', ''10691'', ''10692'', ''10693'', ''10694'', ''10695'', ''10696'', ''10697'', ''10698'', ''10699'', ''107'', ''1070'', ''10700'', ''10701'', ''10702'', ''10703'', ''10704'', ''10705'', ''10706'', ''10707'', ''10708'', ''10709'', ''1071'', ''10710'', ''10711'', ''10712'', ''10713'', ''10714'', ''10715'', ''10716'', ''10717'', ''10718'', ''10719'', ''1072'', ''10720'', ''10721'', ''10722'', ''10723'', ''10724'', ''10725'', ''10726'', ''10727'', ''10728'', ''10729'', ''1073'', ''10730'', ''10731'', ''10732'', ''10733'', ''10734'', ''10735'', ''10736'', ''10737'', ''10738'', ''10739'', ''1074'', ''10740'', ''10741'', ''10742'', ''10743'', ''10744'', ''10745'', ''10746'', ''10747'', ''10748'', ''10749'', ''1075'', ''10750'', ''10751'', ''10752'', ''10753'', ''10754'', ''10755'', ''10756'', ''10757'', ''10758'', ''10759'', ''1076'', ''10760'', ''10761'', ''10762'', ''10763'', ''10764'', ''10765'', ''10766'', ''10767'', ''10768'', ''10769'', ''1077'', ''10770'', ''10771'', ''10772'', ''10773'', ''10774'', ''10775'', ''10776'', ''10777'', ''10778'', ''10779'', ''1078'', ''10780'', ''10781'', ''10782'', ''10783'', ''10784'', ''10785'', ''10786'', ''10787'', ''10788'', ''10789'', ''1079'', ''10790'', ''10791'', ''10792'', ''10793'', ''10794'', ''10795'', ''10796'', ''10797'', ''10798'', ''10799'', ''108'', ''1080'', ''10800'', ''10801'', ''10802'', ''10803'', ''10804'', ''10805'', ''10806'', ''10807'', ''10808'', ''10809'', ''1081'', ''10810'', ''10811'', ''10812'', ''10813'', ''10814'', ''10815'', ''10816'', ''10817'', ''10818'', ''10819'', ''1082'', ''10820'', ''10821'', ''10822'', ''10823'', ''10824'', ''10825'', ''10826'', ''10827'', ''10828'', ''10829'', ''1083'', ''10830'', ''10831'', ''10832'', ''10833'', ''10834'', ''10835'', ''10836'', ''10837'', ''10838'', ''10839'', ''1084'', ''10840'', ''10841'', ''10842'', ''10843'', ''10844'', ''10845'', ''10846'', ''10847'', ''10848'', ''10849'', ''1085'', ''10850'', ''10851'', ''10852'', ''10853'', ''10854'', ''10855'', ''10856'', ''10857'', ''10858'', ''10859'', ''1086'', ''10860'', ''10861'', ''10862'', ''10863'', ''10864'', ''10865'', ''10866'', ''10867'', ''10868'', ''10869'', ''1087'', ''10870'', ''10871'', ''10872'', ''10873'', ''10874'', ''10875'', ''10876'', ''10877'', ''10878'', ''10879'', ''1088'', ''10880'', ''10881'', ''10882'', ''10883'', ''10884'', ''10885'', ''10886'', ''10887'', ''10888'', ''10889'', ''1089'', ''10890'', ''10891'', ''10892'', ''10893'', ''10894'', ''10895'', ''10896'', ''10897'', ''10898'', ''10899'', ''109'', ''1090'', ''10900'', ''10901'', ''10902'', ''10903'', ''10904'', ''10905'', ''10906'', ''10907'', ''10908'', ''10909'', ''1091'', ''10910'', ''10911'', ''10912'', ''10913'', ''10914'', ''10915'', ''10916'', ''10917'', ''10918'', ''10919'', ''1092'', ''10920'', ''10921'', ''10922'', ''10923'', ''10924'', ''10925'', ''10926'', ''10927'', ''10928'', ''10929'', ''1093'', ''10930'', ''10931'', ''10932'', ''10933'', ''10934'', ''10935'', ''10936'', ''10937'', ''10938'', ''10939'', ''1094'', ''10940'', ''10941'', ''10942'', ''10943'', ''10944'', ''10945'', ''10946'', ''10947'', ''10948'', ''10949'', ''1095'', ''10950'', ''10951'', ''10952'', ''10953'', ''10954'', ''10955'', ''10956'', ''10957'', ''10958'', ''10959'', ''1096'', ''10960'', ''10961'', ''10962'', ''10963'', ''10964'', ''10965'', ''10966'', ''10967'', ''10968'', ''10969'', ''1097'', ''10970'', ''10971'', ''10972'', ''10973'', ''10974'', ''10975'', ''10976'', ''10977'', ''10978'', ''10979'', ''1098'', ''10980'', ''10981'', ''10982'', ''10983'', ''10984'', ''10985'', ''10986'', ''10987'', ''10988'', ''10989'', ''1099'', ''10990'', ''10991'', ''10992'', ''10993'', ''10994'', ''10995'', ''10996'', ''10997'', ''10998'', ''10999'', ''11'', ''110'', ''1100'', ''11000'', ''11001'', ''11002'', ''11003'', ''11004'', ''11005'', ''11006'', ''11007'', ''11008'', ''11009'', ''1101'', ''11010'', ''11011'', ''11012'', ''11013'', ''11014'', ''11015'', ''11016'', ''11017'', ''11018'', ''11019'', ''1102'', ''11020'', ''11021'', ''11022'', ''11023'', ''11024'', ''11025'', ''11026'', ''11027'', ''11028'', ''11029'', ''1103'', ''11030'', ''11031'', ''11032'', ''11033'', ''11034'', ''11035'', ''11036'', ''11037'', ''11038'', ''11039'', ''1104'', ''11040'', ''11041'', ''11042'', ''11043'', ''11044'', ''11045'', ''11046'', ''11047'', ''11048'', ''11049'', ''1105'', ''11050'', ''11051'', ''11052'', ''11053'', ''11054'', ''11055'', ''11056'', ''11057'', ''11058'', ''11059'', ''1106'', ''11060'', ''11061'', ''11062'', ''11063'', ''11064'', ''11065'', ''11066'', ''11067'', ''11068'', ''11069'', ''1107'', ''11070'', ''11071'', ''11072'', ''11073'', ''11074'', ''11075'', ''11076'', ''11077'', ''11078'', ''11079'', ''1108'', ''11080'', ''11081'', ''11082'', ''11083'', ''11084'', ''11085'', ''11086'', ''11087'', ''11088'', ''11089'', ''1109'', ''11090'', ''11091'', ''11092'', ''11093'', ''11094'', ''11095'', ''11096'', ''11097'', ''11098'', ''11099'', ''111'', ''1110'', ''11100'', ''11101'', ''11102'', ''11103'', ''11104'', ''11105'', ''11106'', ''11107'', ''11108'', ''11109'', ''1111'', ''11110'', ''11111'', ''11112'', ''11113'', ''11114'', ''11115'', ''11116'', ''11117'', ''11118'', ''11119'', ''1112'', ''11120'', ''11121'', ''11122'', ''11123'', ''11124'', ''11125'', ''11126'', ''11127'', ''11128'', ''11129'', ''1113'', ''11130'', ''11131'', ''11132'', ''11133'', ''11134'', ''11135'', ''11136'', ''11137'', ''11138'', ''11139'', ''1114'', ''11140'', ''11141'', ''11142'', ''11143'', ''11144'', ''11145'', ''11146'', ''11147'', ''11148'', ''11149'', ''1115'', ''11150'', ''11151'', ''11152'', ''11153'', ''11154'', ''11155'', ''11156'', ''11157'', ''11158'', ''11159'', ''1116'', ''11160'', ''11161'', ''11162'', ''11163'', ''11164'', ''11165'', ''11166'', ''11167'', ''11168'', ''11169'', ''1117'', ''11170'', ''11171'', ''11172'', ''11173'', ''11174'', ''11175'', ''11176'', ''11177'', ''11178'', ''11179'', ''1118'', ''11180'', ''11181'', ''11182'', ''11183'', ''11184'', ''11185'', ''11186'', ''11187'', ''11188'', ''11189'', ''1119'', ''11190'', ''11191'', ''11192'', ''11193'', ''11194'', ''11195'', ''11196'', ''11197'', ''11198'', ''11199'', ''112'', ''1120'', ''11200'', ''11201'', ''11202'', ''11203'', ''11204'', ''11205'', ''11206'', ''11207'', ''11208'', ''11209'', ''1121'', ''11210'', ''11211'', ''11212'', ''11213'', ''11214'', ''11215'', ''11216'', ''11217'', ''11218'', ''11219'', ''1122'', ''11220'', ''11221'', ''11222'', ''11223'', ''11224'', ''11225'', ''11226'', ''11227'', ''11228'', ''11229'', ''1123'', ''11230'', ''11231'', ''11232'', ''11233'', ''11234'', ''11235'', ''11236'', ''11237'', ''11238'', ''11239'', ''1124'', ''11240'', ''11241'', ''11242'', ''11243'', ''11244'', ''11245'', ''11246'', ''11247'', ''11248'', ''11249'', ''1125'', ''11250'', ''11251'', ''11252'', ''11253'', ''11254'', ''11255'', ''11256'', ''11257'', ''11258'', ''11259'', ''1126'', ''11260'', ''11261'', ''11262'', ''11263'', ''11264'', ''11265'', ''11266'', ''11267'', ''11268'', ''11269'', ''1127'', ''11270'', ''11271'', ''11272'', ''11273'', ''11274'', ''11275'', ''11276'', ''11277'', ''11278'', ''11279'', ''1128'', ''11280'', ''11281'', ''11282'', ''11283'', ''11284'', ''11285'', ''11286'', ''11287'', ''11288'', ''11289'', ''1129'', ''11290'', ''11291'', ''11292'', ''11293'', ''11294'', ''11295'', ''11296'', ''11297'', ''11298'', ''11299'', ''113'', ''1130'', ''11300'', ''11301'', ''11302'', ''11303'', ''11304'', ''11305'', ''11306'', ''11307'', ''11308'', ''11309'', ''1131'', ''11310'', ''11311'', ''11312'', ''11313'', ''11314'', ''11315'', ''11316'', ''11317'', ''11318'', ''11319'', ''1132'', ''11320'', ''11321'', ''11322'', ''11323'', ''11324'', ''11325'', ''11326'', ''11327'', ''11328'', ''11329'', ''1133'', ''11330'', ''11331'', ''11332'', ''11333'', ''11334'', ''11335'', ''11336'', ''11337'', ''11338'', ''11339'', ''1134'', ''11340'', ''11341'', ''11342'', ''11343'', ''11344'', ''11345'', ''11346'', ''11347'', ''11348'', ''11349'', ''1135'', ''11350'', ''11351'', ''11352'', ''11353'', ''11354'', ''11355'', ''11356'', ''11357'', ''11358'', ''11359'', ''1136'', ''11360'', ''11361'', ''11362'', ''11363'', ''11364'', ''11365'', ''11366'', ''11367'', ''11368'', ''11369'', ''1137'', ''11370'', ''11371'', ''11372'', ''11373'', ''11374'', ''11375'', ''11376'', ''11377'', ''11378'', ''11379'', ''1138'', ''11380'', ''11381'', ''11382'', ''11383'', ''11384'', ''11385'', ''11386'', ''11387'', ''11388'', ''11389'', ''1139'', ''11390'', ''11391'', ''11392'', ''11393'', ''11394'', ''11395'', ''11396'', ''11397'', ''11398'', ''11399'', ''114'', ''1140'', ''11400'', ''11401'', ''11402'', ''11403'', ''11404'', ''11405'', ''11406'', ''11407'', ''11408'', ''11409'', ''1141'', ''11410'', ''11411'', ''11412'', ''11413'', ''11414'', ''11415'', ''11416'', ''11417'', ''11418'', ''11419'', ''1142'', ''11420'', ''11421'', ''11422'', ''11423'', ''11424'', ''11425'', ''11426'', ''11427'', ''11428'', ''11429'', ''1143'', ''11430'', ''11431'', ''11432'', ''11433'', ''11434'', ''11435'', ''11436'', ''11437'', ''11438'', ''11439'', ''1144'', ''11440'', ''11441'', ''11442'', ''11443'', ''11444'', ''11445'', ''11446'', ''11447'', ''11448'', ''11449'', ''1145'', ''11450'', ''11451'', ''11452'', ''11453'', ''11454'', ''11455'', ''11456'', ''11457'', ''11458'', ''11459'', ''1146'', ''11460'', ''11461'', ''11462'', ''11463'', ''11464'', ''11465'', ''11466'', ''11467'', ''11468'', ''11469'', ''1147'', ''11470'', ''11471'', ''11472'', ''11473'', ''11474'', ''11475'', ''11476'', ''11477'', ''11478'', ''11479'', ''1148'', ''11480'', ''11481'', ''11482'', ''11483'', ''11484'', ''11485'', ''11486'', ''11487'', ''11488'', ''11489'', ''1149'', ''11490'', ''11491'', ''11492'', ''11493'', ''11494'', ''11495'', ''11496'', ''11497'', ''11498'', ''11499'', ''115'', ''1150'', ''11500'', ''11501'', ''11502'', ''11503'', ''11504'', ''11505'', ''11506'', ''11507'', ''11508'', ''11509'', ''1151'', ''11510'', ''11511'', ''11512'', ''11513'', ''11514'', ''11515'', ''11516'', ''11517'', ''11518'', ''11519'', ''1152'', ''11520'', ''11521'', ''11522'', ''11523'', ''11524'', ''11525'', ''11526'', ''11527'', ''11528'', ''11529'', ''1153'', ''11530'', ''11531'', ''11532'', ''11533'', ''11534'', ''11535'', ''11536'', ''11537'', ''11538'', ''11539'', ''1154'', ''11540'', ''11541'', ''11542'', ''11543'', ''11544'', ''11545'', ''11546'', ''11547'', ''11548'', ''11549'', ''1155'', ''11550'', ''11551'', ''11552'', ''11553'', ''11554'', ''11555'', ''11556'', ''11557'', ''11558'', ''11559'', ''1156'', ''11560'', ''11561'', ''11562'', ''11563'', ''11564'', ''11565'', ''11566'', ''11567'', ''11568'', ''11569'', ''1157'', ''11570'', ''11571'', ''11572'', ''11573'', ''11574'', ''11575'', ''11576'', ''11577'', ''11578'', ''11579'', ''1158'', ''11580'', ''11581'', ''11582'', ''11583'', ''11584'', ''11585'', ''11586'', ''11587'', ''11588'', ''11589'', ''1159'', ''11590'', ''11591'', ''11592'', ''11593'', ''11594'', ''11595'', ''11596'', ''11597'', ''11598'', ''11599'', ''116'', ''1160'', ''11600'', ''11601'', ''11602'', ''11603'', ''11604'', ''11605'', ''11606'', ''11607'', ''11608'', ''11609'', ''1161'', ''11610'', ''11611'', ''11612'', ''11613'', ''11614'', ''11615'', ''11616'', ''11617'', ''11618'', ''11619'', ''1162'', ''11620'', ''11621'', ''11622'', ''11623'', ''11624'', ''11625'', ''11626'', ''11627'', ''11628'', ''11629'', ''1163'', ''11630'', ''11631'', ''11632'', ''11633'', ''11634'', ''11635'', ''11636'', ''11637'', ''11638'', ''11639'', ''1164'', ''11640'', ''11641'', ''11642'', ''11643'', ''11644'', ''11645'', ''11646'', ''11647'', ''11648'', ''11649'', ''1165'', ''11650'', ''11651'', ''11652'', ''11653'', ''11654'', ''11655'', ''11656'', ''11657'', ''11658'', ''11659'', ''1166'', ''11660'', ''11661'', ''11662'', ''11663'', ''11664'', ''11665'', ''11666'', ''11667'', ''11668'', ''11669'', ''1167'', ''11670'', ''11671'', ''11672'', ''11673'', ''11674'', ''11675'', ''11676'', ''11677'', ''11678'', ''11679'', ''1168'', ''11680'', ''11681'', ''11682'', ''11683'', ''11684'', ''11685'', ''11686'', ''11687'', ''11688'', ''11689'', ''1169'', ''11690'', ''11691'', ''11692'', ''11693'', ''11694'', ''11695'', ''11696'', ''11697'', ''11698'', ''11699'', ''117'', ''1170'', ''11700'', ''11701'', ''11702'', ''11703'', ''11704'', ''11705'', ''11706'', ''11707'', ''11708'', ''11709'', ''1171'', ''11710'', ''11711'', ''11712'', ''11713'', ''11714'', ''11715'', ''11716'', ''11717'', ''11718'', ''11719'', ''1172'', ''11720'', ''11721'', ''11722'', ''11723'', ''11724'', ''11725'', ''11726'', ''11727'', ''11728'', ''11729'', ''1173'', ''11730'', ''11731'', ''11732'', ''11733'', ''11734'', ''11735'', ''11736'', ''11737'', ''11738'', ''11739'', ''1174'', ''11740'', ''11741'', ''11742'', ''11743'', ''11744'', ''11745'', ''11746'', ''11747'', ''11748'', ''11749'', ''1175'', ''11750'', ''11751'', ''11752'', ''11753'', ''11754'', ''11755'', ''11756'', ''11757'', ''11758'', ''11759'', ''1176'', ''11760'', ''11761'', ''11762'', ''11763'', ''11764'', ''11765'', ''11766'', ''11767'', ''11768'', ''11769'', ''1177'', ''11770'', ''11771'', ''11772'', ''11773'', ''11774'', ''11775'', ''11776'', ''11777'', ''11778'', ''11779'', ''1178'', ''11780'', ''11781'', ''11782'', ''11783'', ''11784'', ''11785'', ''11786'', ''11787'', ''11788'', ''11789'', ''1179'', ''11790'', ''11791'', ''11792'', ''11793'', ''11794'', ''11795'', ''11796'', ''11797'', ''11798'', ''11799'', ''118'', ''1180'', ''11800'', ''11801'', ''11802'', ''11803'', ''11804'', ''11805'', ''11806'', ''11807'', ''11808'', ''11809'', ''1181'', ''11810'', ''11811'', ''11812'', ''11813'', ''11814'', ''11815'', ''11816'', ''11817'', ''11818'', ''11819'', ''1182'', ''11820'', ''11821'', ''11822'', ''11823'', ''11824'', ''11825'', ''11826'', ''11827'', ''11828'', ''11829'', ''1183'', ''11830'', ''11831'', ''11832'', ''11833'', ''11834'', ''11835'', ''11836'', ''11837'', ''11838'', ''11839'', ''1184'', ''11840'', ''11841'', ''11842'', ''11843'', ''11844'', ''11845'', ''11846'', ''11847'', ''11848'', ''11849'', ''1185'', ''11850'', ''11851'', ''11852'', ''11853'', ''11854'', ''11855'', ''11856'', ''11857'', ''11858'', ''11859'', ''1186'', ''11860'', ''11861'', ''11862'', ''11863'', ''11864'', ''11865'', ''11866'', ''11867'', ''11868'', ''11869'', ''1187'', ''11870'', ''11871'', ''11872'', ''11873'', ''11874'', ''11875'', ''11876'', ''11877'', ''11878'', ''11879'', ''1188'', ''11880'', ''11881'', ''11882'', ''11883'', ''11884'', ''11885'', ''11886'', ''11887'', ''11888'', ''11889'', ''1189'', ''11890'', ''11891'', ''11892'', ''11893'', ''11894'', ''11895'', ''11896'', ''11897'', ''11898'', ''11899'', ''119'', ''1190'', ''11900'', ''11901'', ''11902'', ''11903'', ''11904'', ''11905'', ''11906'', ''11907'', ''11908'', ''11909'', ''1191'', ''11910'', ''11911'', ''11912'', ''11913'', ''11914'', ''11915'', ''11916'', ''11917'', ''11918'', ''11919'', ''1192'', ''11920'', ''11921'', ''11922'', ''11923'', ''11924'', ''11925'', ''11926'', ''11927'', ''11928'', ''11929'', ''1193'', ''11930'', ''11931'', ''11932'', ''11933'', ''11934'', ''11935'', ''11936'', ''11937'', ''11938'', ''11939'', ''1194'', ''11940'', ''11941'', ''11942'', ''11943'', ''11944'', ''11945'', ''11946'', ''11947'', ''11948'', ''11949'', ''1195'', ''11950'', ''11951'', ''11952'', ''11953'', ''11954'', ''11955'', ''11956'', ''11957'', ''11958'', ''11959'', ''1196'', ''11960'', ''11961'', ''11962'', ''11963'', ''11964'', ''11965'', ''11966'', ''11967'', ''11968'', ''11969'', ''1197'', ''11970'', ''11971'', ''11972'', ''11973'', ''11974'', ''11975'', ''11976'', ''11977'', ''11978'', ''11979'', ''1198'', ''11980'', ''11981'', ''11982'', ''11983'', ''11984'', ''11985'', ''11986'', ''11987'', ''11988'', ''11989'', ''1199'', ''11990'', ''11991'', ''11992'', ''11993'', ''11994'', ''11995'', ''11996'', ''11997'', ''11998'', ''11999'', ''12'', ''120'', ''1200'', ''12000'', ''12001'', ''12002'', ''12003'', ''12004'', ''12005'', ''12006'', ''12007'', ''12008'', ''12009'', ''1201'', ''12010'', ''12011'', ''12012'', ''12013'', ''12014'', ''12015'', ''12016'', ''12017'', ''12018'', ''12019'', ''1202'', ''12020'', ''12021'', ''12022'', ''12023'', ''12024'', ''12025'', ''12026'', ''12027'', ''12028'', ''12029'', ''1203'', ''12030'', ''12031'', ''12032'', ''12033'', ''12034'', ''12035'', ''12036'', ''12037'', ''12038'', ''12039'', ''1204'', ''12040'', ''12041'', ''12042'', ''12043'', ''12044'', ''12045'', ''12046'', ''12047'', ''12048'', ''12049'', ''1205'', ''12050'', ''12051'', ''12052'', ''12053'', ''12054'', ''12055'', ''12056'', ''12057'', ''12058'', ''12059'', ''1206'', ''12060'', ''12061'', ''12062'', ''12063'', ''12064'', ''12065'', ''12066'', ''12067'', ''12068'', ''12069'', ''1207'', ''12070'', ''12071'', ''12072'', ''12073'', ''12074'', ''12075'', ''12076'', ''12077'', ''12078'', ''12079'', ''1208'', ''12080'', ''12081'', ''12082'', ''12083'', ''12084'', ''12085'', ''12086'', ''12087'', ''12088'', ''12089'', ''1209'', ''12090'', ''12091'', ''12092'', ''12093'', ''12094'', ''12095'', ''12096'', ''12097'', ''12098'', ''12099'', ''121'', ''1210'', ''12100'', ''12101'', ''12102'', ''12103'', ''12104'', ''12105'', ''12106'', ''12107'', ''12108'', ''12109'', ''1211'', ''12110'', ''12111'', ''12112'', ''12113'', ''12114'', ''12115'', ''12116'', ''12117'', ''12118'', ''12119'', ''1212'', ''12120'', ''12121'', ''12122'', ''12123'', ''12124'', ''12125'', ''12126'', ''12127'', ''12128'', ''12129'', ''1213'', ''12130'', ''12131'', ''12132'', ''12133'', ''12134'', ''12135'', ''12136'', ''12137'', ''12138'', ''12139'', ''1214'', ''12140'', ''12141'', ''12142'', ''12143'', ''12144'', ''12145'', ''12146'', ''12147'', ''12148'', ''12149'', ''1215'', ''12150'', ''12151'', ''12152'', ''12153'', ''12154'', ''12155'', ''12156'', ''12157'', ''12158'', ''12159'', ''1216'', ''12160'', ''12161'', ''12162'', ''12163'', ''12164'', ''12165'', ''12166'', ''12167'', ''12168'', ''12169'', ''1217'', ''12170'', ''12171'', ''12172'', ''12173'', ''12174'', ''12175'', ''12176'', ''12177'', ''12178'', ''12179'', ''1218'', ''12180'', ''12181'', ''12182'', ''12183'', ''12184'', ''12185'', ''12186'', ''12187'', ''12188'', ''12189'', ''1219'', ''12190'', ''12191'', ''12192'', ''12193'', ''12194'', ''12195'', ''12196'', ''12197'', ''12198'', ''12199'', ''122'', ''1220'', ''12200'', ''12201'', ''12202'', ''12203'', ''12204'', ''12205'', ''12206'', ''12207'', ''12208'', ''12209'', ''1221'', ''12210'', ''12211'', ''12212'', ''12213'', ''12214'', ''12215'', ''12216'', ''12217'', ''12218'', ''12219'', ''1222'', ''12220'', ''12221'', ''12222'', ''12223'', ''12224'', ''12225'', ''12226'', ''12227'', ''12228'', ''12229'', ''1223'', ''12230'', ''12231'', ''12232'', ''12233'', ''12234'', ''12235'', ''12236'', ''12237'', ''12238'', ''12239'', ''1224'', ''12240'', ''12241'', ''12242'', ''12243'', ''12244'', ''12245'', ''12246'', ''12247'', ''12248'', ''12249'', ''1225'', ''12250'', ''12251'', ''12252'', ''12253'', ''12254'', ''12255'', ''12256'', ''12257'', ''12258'', ''12259'', ''1226'', ''12260'', ''12261'', ''12262'', ''12263'', ''12264'', ''12265'', ''12266'', ''12267'', ''12268'', ''12269'', ''1227'', ''12270'', ''12271'', ''12272'', ''12273'', ''12274'', ''12275'', ''12276'', ''12277'', ''12278'', ''12279'', ''1228'', ''12280'', ''12281'', ''12282'', ''12283'', ''12284'', ''12285'', ''12286'', ''12287'', ''12288'', ''12289'', ''1229'', ''12290'', ''12291'', ''12292'', ''12293'', ''12294'', ''12295'', ''12296'', ''12297'', ''12298'', ''12299'', ''123'', ''1230'', ''12300'', ''12301'', ''12302'', ''12303'', ''12304'', ''12305'', ''12306'', ''12307'', ''12308'', ''12309'', ''1231'', ''12310'', ''12311'', ''12312'', ''12313'', ''12314'', ''12315'', ''12316'', ''12317'', ''12318'', ''12319'', ''1232'', ''12320'', ''12321'', ''12322'', ''12323'', ''12324'', ''12325'', ''12326'', ''12327'', ''12328'', ''12329'', ''1233'', ''12330'', ''12331'', ''12332'', ''12333'', ''12334'', ''12335'', ''12336'', ''12337'', ''12338'', ''12339'', ''1234'', ''12340'', ''12341'', ''12342'', ''12343'', ''12344'', ''12345'', ''12346'', ''12347'', ''12348'', ''12349'', ''1235'', ''12350'', ''12351'', ''12352'', ''12353'', ''12354'', ''12355'', ''12356'', ''12357'', ''12358'', ''12359'', ''1236'', ''12360'', ''12361'', ''12362'', ''12363'', ''12364'', ''12365'', ''12366'', ''12367'', ''12368'', ''12369'', ''1237'', ''12370'', ''12371'', ''12372'', ''12373'', ''12374'', ''12375'', ''12376'', ''12377'', ''12378'', ''12379'', ''1238'', ''12380'', ''12381'', ''12382'', ''12383'', ''12384'', ''12385'', ''12386'', ''12387'', ''12388'', ''12389'', ''1239'', ''12390'', ''12391'', ''12392'', ''12393'', ''12394'', ''12395'', ''12396'', ''12397'', ''12398'', ''12399'', ''124'', ''1240'', ''12400'', ''12401'', ''12402'', ''12403'', ''12404'', ''12405'', ''12406'', ''12407'', ''12408'', ''12409'', ''1241'', ''12410'', ''12411'', ''12412'', ''12413'', ''12414'', ''12415'', ''12416'', ''12417'', ''12418'', ''12419'', ''1242'', ''12420'', ''12421'', ''12422'', ''12423'', ''12424'', ''12425'', ''12426'', ''12427'', ''12428'', ''12429'', ''1243'', ''12430'', ''12431'', ''12432'', ''12433'', ''12434'', ''12435'', ''12436'', ''12437'', ''12438'', ''12439'', ''1244'', ''12440'', ''12441'', ''12442'', ''12443'', ''12444'', ''12445'', ''12446'', ''12447'', ''12448'', ''12449'', ''1245'', ''12450'', ''12451'', ''12452'', ''12453'', ''12454'', ''12455'', ''12456'', ''12457'', ''12458'', ''12459'', ''1246'', ''12460'', ''12461'', ''12462'', ''12463'', ''12464'', ''12465'', ''12466'', ''12467'', ''12468'', ''12469'', ''1247'', ''12470'', ''12471'', ''12472'', ''12473'', ''12474'', ''12475'', ''12476'', ''12477'', ''12478'', ''12479'', ''1248'', ''12480'', ''12481'', ''12482'', ''12483'', ''12484'', ''12485'', ''12486'', ''12487'', ''12488'', ''12489'', ''1249'', ''12490'', ''12491'', ''12492'', ''12493'', ''12494'', ''12495'', ''12496'', ''12497'', ''12498'', ''12499'', ''125'', ''1250'', ''12500'', ''12501'', ''12502'', ''12503'', ''12504'', ''12505'', ''12506'', ''12507'', ''12508'', ''12509'', ''1251'', ''12510'', ''12511'', ''12512'', ''12513'', ''12514'', ''12515'', ''12516'', ''12517'', ''12518'', ''12519'', ''1252'', ''12520'', ''12521'', ''12522'', ''12523'', ''12524'', ''12525'', ''12526'', ''12527'', ''12528'', ''12529'', ''1253'', ''12530'', ''12531'', ''12532'', ''12533'', ''12534'', ''12535'', ''12536'', ''12537'', ''12538'', ''12539'', ''1254'', ''12540'', ''12541'', ''12542'', ''12543'', ''12544'', ''12545'', ''12546'', ''12547'', ''12548'', ''12549'', ''1255'', ''12550'', ''12551'', ''12552'', ''12553'', ''12554'', ''12555'', ''12556'', ''12557'', ''12558'', ''12559'', ''1256'', ''12560'', ''12561'', ''12562'', ''12563'', ''12564'', ''12565'', ''12566'', ''12567'', ''12568'', ''12569'', ''1257'', ''12570'', ''12571'', ''12572'', ''12573'', ''12574'', ''12575'', ''12576'', ''12577'', ''12578'', ''12579'', ''1258'', ''12580'', ''12581'', ''12582'', ''12583'', ''12584'', ''12585'', ''12586'', ''12587'', ''12588'', ''12589'', ''1259'', ''12590'', ''12591'', ''12592'', ''12593'', ''12594'', ''12595'', ''12596'', ''12597'', ''12598'', ''12599'', ''126'', ''1260'', ''12600'', ''12601'', ''12602'', ''12603'', ''12604'', ''12605'', ''12606'', ''12607'', ''12608'', ''12609'', ''1261'', ''12610'', ''12611'', ''12612'', ''12613'', ''12614'', ''12615'', ''12616'', ''12617'', ''12618'', ''12619'', ''1262'', ''12620'', ''12621'', ''12622'', ''12623'', ''12624'', ''12625'', ''12626'', ''12627'', ''12628'', ''12629'', ''1263'', ''12630'', ''12631'', ''12632'', ''12633'', ''12634'', ''12635'', ''12636'', ''12637'', ''12638'', ''12639'', ''1264'', ''12640'', ''12641'', ''12642'', ''12643'', ''12644'', ''12645'', ''12646'', ''12647'', ''12648'', ''12649'', ''1265'', ''12650'', ''12651'', ''12652'', ''12653'', ''12654'', ''12655'', ''12656'', ''12657'', ''12658'', ''12659'', ''1266'', ''12660'', ''12661'', ''12662'', ''12663'', ''12664'', ''12665'', ''12666'', ''12667'', ''12668'', ''12669'', ''1267'', ''12670'', ''12671'', ''12672'', ''12673'', ''12674'', ''12675'', ''12676'', ''12677'', ''12678'', ''12679'', ''1268'', ''12680'', ''12681'', ''12682'', ''12683'', ''12684'', ''12685'', ''12686'', ''12687'', ''12688'', ''12689'', ''1269'', ''12690'', ''12691'', ''12692'', ''12693'', ''12694'', ''12695'', ''12696'', ''12697'', ''12698'', ''12699'', ''127'', ''1270'', ''12700'', ''12701'', ''12702'', ''12703'', ''12704'', ''12705'', ''12706'', ''12707'', ''12708'', ''12709'', ''1271'', ''12710'', ''12711'', ''12712'', ''12713'', ''12714'', ''12715'', ''12716'', ''12717'', ''12718'', ''12719'', ''1272'', ''12720'', ''12721'', ''12722'', ''12723'', ''12724'', ''12725'', ''12726'', ''12727'', ''12728'', ''12729'', ''1273'', ''12730'', ''12731'', ''12732'', ''12733'', ''12734'', ''12735'', ''12736'', ''12737'', ''12738'', ''12739'', ''1274'', ''12740'', ''12741'', ''12742'', ''12743'', ''12744'', ''12745'', ''12746'', ''12747'', ''12748'', ''12749'', ''1275'', ''12750'', ''12751'', ''12752'', ''12753'', ''12754'', ''12755'', ''12756'', ''12757'', ''12758'', ''12759'', ''1276'', ''12760'', ''12761'', ''12762'', ''12763'', ''12764'', ''12765'', ''12766'', ''12767'', ''12768'', ''12769'', ''1277'', ''12770'', ''12771'', ''12772'', ''12773'', ''12774'', ''12775'', ''12776'', ''12777'', ''12778'', ''12779'', ''1278'', ''12780'', ''12781'', ''12782'', ''12783'', ''12784'', ''12785'', ''12786'', ''12787'', ''12788'', ''12789'', ''1279'', ''12790'', ''12791'', ''12792'', ''12793'', ''12794'', ''12795'', ''12796'', ''12797'', ''12798'', ''12799'', ''128'', ''1280'', ''12800'', ''12801'', ''12802'', ''12803'', ''12804'', ''12805'', ''12806'', ''12807'', ''12808'', ''12809'', ''1281'', ''12810'', ''12811'', ''12812'', ''12813'', ''12814'', ''12815'', ''12816'', ''12817'', ''12818'', ''12819'', ''1282'', ''12820'', ''12821'', ''12822'', ''12823'', ''12824'', ''12825'', ''12826'', ''12827'', ''12828'', ''12829'', ''1283'', ''12830'', ''12831'', ''12832'', ''12833'', ''12834'', ''12835'', ''12836'', ''12837'', ''12838'', ''12839'', ''1284'', ''12840'', ''12841'', ''12842'', ''12843'', ''12844'', ''12845'', ''12846'', ''12847'', ''12848'', ''12849'', ''1285'', ''12850'', ''12851'', ''12852'', ''12853'', ''12854'', ''12855'', ''12856'', ''12857'', ''12858'', ''12859'', ''1286'', ''12860'', ''12861'', ''12862'', ''12863'', ''12864'', ''12865'', ''12866'', ''12867'', ''12868'', ''12869'', ''1287'', ''12870'', ''12871'', ''12872'', ''12873'', ''12874'', ''12875'', ''12876'', ''12877'', ''12878'', ''12879'', ''1288'', ''12880'', ''12881'', ''12882'', ''12883'', ''12884'', ''12885'', ''12886'', ''12887'', ''12888'', ''12889'', ''1289'', ''12890'', ''12891'', ''12892'', ''12893'', ''12894'', ''12895'', ''12896'', ''12897'', ''12898'', ''12899'', ''129'', ''1290'', ''12900'', ''12901'', ''12902'', ''12903'', ''12904'', ''12905'', ''12906'', ''12907'', ''12908'', ''12909'', ''1291'', ''12910'', ''12911'', ''12912'', ''12913'', ''12914'', ''12915'', ''12916'', ''12917'', ''12918'', ''12919'', ''1292'', ''12920'', ''12921'', ''12922'', ''12923'', ''12924'', ''12925'', ''12926'', ''12927'', ''12928'', ''12929'', ''1293'', ''12930'', ''12931'', ''12932'', ''12933'', ''12934'', ''12935'', ''12936'', ''12937'', ''12938'', ''12939'', ''1294'', ''12940'', ''12941'', ''12942'', ''12943'', ''12944'', ''12945'', ''12946'', ''12947'', ''12948'', ''12949'', ''1295'', ''12950'', ''12951'', ''12952'', ''12953'', ''12954'', ''12955'', ''12956'', ''12957'', ''12958'', ''12959'', ''1296'', ''12960'', ''12961'', ''12962'', ''12963'', ''12964'', ''12965'', ''12966'', ''12967'', ''12968'', ''12969'', ''1297'', ''12970'', ''12971'', ''12972'', ''12973'', ''12974'', ''12975'', ''12976'', ''12977'', ''12978'', ''12979'', ''1298'', ''12980'', ''12981'', ''12982'', ''12983'', ''12984'', ''12985'', ''12986'', ''12987'', ''12988'', ''12989'', ''1299'', ''12990'', ''12991'', ''12992'', ''12993'', ''12994'', ''12995'', ''12996'', ''12997'', ''12998'', ''12999'', ''13'', ''130'', ''1300'', ''13000'', ''13001'', ''13002'', ''13003'', ''13004'', ''13005'', ''13006'', ''13007'', ''13008'', ''13009'', ''1301'', ''13010'', ''13011'', ''13012'', ''13013'', ''13014'', ''13015'', ''13016'', ''13017'', ''13018'', ''13019'', ''1302'', ''13020'', ''13021'', ''13022'', ''13023'', ''13024'', ''13025'', ''13026'', ''13027'', ''13028'', ''13029'', ''1303'', ''13030'', ''13031'', ''13032'', ''13033'', ''13034'', ''13035'', ''13036'', ''13037'', ''13038'', ''13039'', ''1304'', ''13040'', ''13041'', ''13042'', ''13043'', ''13044'', ''13045'', ''13046'', ''13047'', ''13048'', ''13049'', ''1305'', ''13050'', ''13051'', ''13052'', ''13053'', ''13054'', ''13055'', ''13056'', ''13057'', ''13058'', ''13059'', ''1306'', ''13060'', ''13061'', ''13062'', ''13063'', ''13064'', ''13065'', ''13066'', ''13067'', ''13068'', ''13069'', ''1307'', ''13070'', ''13071'', ''13072'', ''13073'', ''13074'', ''13075'', ''13076'', ''13077'', ''13078'', ''13079'', ''1308'', ''13080'', ''13081'', ''13082'', ''13083'', ''13084'', ''13085'', ''13086'', ''13087'', ''13088'', ''13089'', ''1309'', ''13090'', ''13091'', ''13092'', ''13093'', ''13094'', ''13095'', ''13096'', ''13097'', ''13098'', ''13099'', ''131'', ''1310'', ''13100'', ''13101'', ''13102'', ''13103'', ''13104'', ''13105'', ''13106'', ''13107'', ''13108'', ''13109'', ''1311'', ''13110'', ''13111'', ''13112'', ''13113'', ''13114'', ''13115'', ''13116'', ''13117'', ''13118'', ''13119'', ''1312'', ''13120'', ''13121'', ''13122'', ''13123'', ''13124'', ''13125'', ''13126'', ''13127'', ''13128'', ''13129'', ''1313'', ''13130'', ''13131'', ''13132'', ''13133'', ''13134'', ''13135'', ''13136'', ''13137'', ''13138'', ''13139'', ''1314'', ''13140'', ''13141'', ''13142'', ''13143'', ''13144'', ''13145'', ''13146'', ''13147'', ''13148'', ''13149'', ''1315'', ''13150'', ''13151'', ''13152'', ''13153'', ''13154'', ''13155'', ''13156'', ''13157'', ''13158'', ''13159'', ''1316'', ''13160'', ''13161'', ''13162'', ''13163'', ''13164'', ''13165'', ''13166'', ''13167'', ''13168'', ''13169'', ''1317'', ''13170'', ''13171'', ''13172'', ''13173'', ''13174'', ''13175'', ''13176'', ''13177'', ''13178'', ''13179'', ''1318'', ''13180'', ''13181'', ''13182'', ''13183'', ''13184'', ''13185'', ''13186'', ''13187'', ''13188'', ''13189'', ''1319'', ''13190'', ''13191'', ''13192'', ''13193'', ''13194'', ''13195'', ''13196'', ''13197'', ''13198'', ''13199'', ''132'', ''1320'', ''13200'', ''13201'', ''13202'', ''13203'', ''13204'', ''13205'', ''13206'', ''13207'', ''13208'', ''13209'', ''1321'', ''13210'', ''13211'', ''13212'', ''13213'', ''13214'', ''13215'', ''13216'', ''13217'', ''13218'', ''13219'', ''1322'', ''13220'', ''13221'', ''13222'', ''13223'', ''13224'', ''13225'', ''13226'', ''13227'', ''13228'', ''13229'', ''1323'', ''13230'', ''13231'', ''13232'', ''13233'', ''13234'', ''13235'', ''13236'', ''13237'', ''13238'', ''13239'', ''1324'', ''13240'', ''13241'', ''13242'', ''13243'', ''13244'', ''13245'', ''13246'', ''13247'', ''13248'', ''13249'', ''1325'', ''13250'', ''13251'', ''13252'', ''13253'', ''13254'', ''13255'', ''13256'', ''13257'', ''13258'', ''13259'', ''1326'', ''13260'', ''13261'', ''13262'', ''13263'', ''13264'', ''13265'', ''13266'', ''13267'', ''13268'', ''13269'', ''1327'', ''13270'', ''13271'', ''13272'', ''13273'', ''13274'', ''13275'', ''13276'', ''13277'', ''13278'', ''13279'', ''1328'', ''13280'', ''13281'', ''13282'', ''13283'', ''13284'', ''13285'', ''13286'', ''13287'', ''13288'', ''13289'', ''1329'', ''13290'', ''13291'', ''13292'', ''13293'', ''13294'', ''13295'', ''13296'', ''13297'', ''13298'', ''13299'', ''133'', ''1330'', ''13300'', ''13301'', ''13302'', ''13303'', ''13304'', ''13305'', ''13306'', ''13307'', ''13308'', ''13309'', ''1331'', ''13310'', ''13311'', ''13312'', ''13313'', ''13314'', ''13'

# Row 3 (surname) updates
$ws.Range("P3").Value = 2
$ws.Range("R3").Value = 654898
$ws.Range("S3").Value = 21926
$ws.Range("T3").Value = '[''0'', ''100 Acres Lots'', ''8th range of Grenville'', ''ABBA'', ''ABBE'', ''ABBEE'', ''ABBOT'', ''ABBOTT'', ''ABBY'', ''ABEL'', ''ABER'', ''ABERCROMBIE'', ''ABERNATHER'', ''ABERNETHY'', ''ABRAHAM'', ''ABRAHAM/COLE'', ''ABRAHAMS'', ''ABRAHAMS/COLE'', ''ABRAM'', ''ABRAMS'', ''ABRAMS/SNYDER'', ''ABRAY'', ''ABURNATHY'', ''ACCARMAN'', ''ACHESON'', ''ACHISON'', ''ACKER'', ''ACKLAND'', ''ACKLEY'', ''ACKMAN'', ''ACOR'', ''ACRE'', ''ACTON'', ''ADAIR'', ''ADAIR/CRAWLEY'', ''ADAM'', ''ADAM/LARAMEY'', ''ADAMS'', ''ADAMS/ADNANS'', ''ADAMS/CLIFFORD'', ''ADAMS/COONS'', ''ADAMSON'', ''ADDISON'', ''ADKINS'', ''ADKINSON'', ''ADSIT/HAIRLAND'', ''ADY'', ''ADZIT'', ''AGER'', ''AGLER'', ''AGLOR/ACHLER'', ''AGNEW'', ''AHIER/ALICER'', ''AHSMAN'', ''AHSTON'', ''AICKMAN'', ''AIKEN'', ''AIKIN'', ''AIKINS'', ''AIKMAN'', ''AILEY'', ''AINLAY'', ''AINSE'', ''AINSE/MONTURE'', ''AINSLEY'', ''AINSLIE'', ''AINSWORTH'', ''AIRD'', ''AIREY'', ''AIREYS'', ''AIRHART'', ''AIRHEART'', ''AIRTH'', ''AISSANCE'', ''AITCHISON'', ''AITHENS'', ''AITKEN'', ''AITKENS'', ''AITKIN'', ''AITKINS'', ''AKE'', ''AKERLY'', ''AKEY'', ''AKIN'', ''AKINS'', ''AKROYD'', ''ALAN'', ''ALBERSON'', ''ALBERT'', ''ALBERTSON'', ''ALBERTSON/POWELL/POWLY'', ''ALBRANT'', ''ALBRIGHT'', ''ALCOCK'', ''ALCOMBRACK'', ''ALCORN'', ''ALDER'', ''ALDERICE'', ''ALDERTON'', ''ALDGEO'', ''ALDRICH'', ''ALDRICK'', ''ALDWORTH'', ''ALEAY'', ''ALEXANDER'', ''ALFORD'', ''ALGAE'', ''ALGEA'', ''ALGEO'', ''ALGER'', ''ALGIER'', ''ALGIO'', ''ALGIRE'', ''ALGOE'', ''ALGUINE'', ''ALGUIRE'', ''ALISON'', ''ALLAIR'', ''ALLAN'', ''ALLARD'', ''ALLAWAY'', ''ALLBRIGHT'', ''ALLCOCK'', ''ALLCOMBRACK'', ''ALLCOTT'', ''ALLCOX'', ''ALLD'', ''ALLEN'', ''ALLEN/COLE'', ''ALLEN/COOPER'', ''ALLEN/CUMMING'', ''ALLEN/MOSHER'', ''ALLENDER'', ''ALLEY'', ''ALLGER'', ''ALLIBONE'', ''ALLICE'', ''ALLIN'', ''ALLIS'', ''ALLISON'', ''ALLPORT'', ''ALLRIGHT'', ''ALLSOP'', ''ALLWOOD'', ''ALMAS'', ''ALMIS'', ''ALMS'', ''ALPHOUGH'', ''ALSTEN'', ''ALSTON'', ''ALT'', ''ALTENBURG'', ''ALTHOUSE'', ''ALTHOUSER'', ''ALTON'', ''ALWARD'', ''ALWAY'', ''ALWOOD'', ''ALYEA'', ''AMABLE'', ''AMAN'', ''AMBLER'', ''AMBRIDGE'', ''AMBROSE'', ''AMELYER'', ''AMERY'', ''AMES'', ''AMEY'', ''AMEY/BAKER'', ''AMIOTTE'', ''AMLEA'', ''AMLIN'', ''AMMERMAN'', ''AMMON'', ''AMOR'', ''AMORY'', ''AMOS'', ''AMPHLETT'', ''AMSTED'', ''AMUS'', ''AMY'', ''AMY/DAVY'', ''ANABELD'', ''ANCASTER'', ''ANCASTER/CROOKS'', ''ANCASTER/CRYSLER'', ''ANDERSON'', ''ANDERSON/CLARKE'', ''ANDERSON/CRYDERMAN'', ''ANDERSON/CUMMINGS'', ''ANDERSON/CURTIS'', "ANDERSON/M''BEAN", ''ANDERSON/McGRATH'', ''ANDERTON'', ''ANDRASS'', ''ANDRES'', ''ANDRESS'', ''ANDREW'', ''ANDREWS'', ''ANDRICK'', ''ANDRUS'', ''ANDRUSS'', ''ANGELL'', ''ANGER'', ''ANGLEMIRE'', ''ANGUISH'', ''ANNABLE'', ''ANNANCE'', ''ANNAND'', ''ANNET'', ''ANNETT'', ''ANNIS'', ''ANSLEY'', ''ANSON'', ''ANTAILLA'', ''ANTALIA'', ''ANTALLIA'', ''ANTAYA'', ''ANTHON'', ''ANTHONY'', ''APPELBEE'', ''APPELLEE'', ''APPLEBEE'', ''APPLEBY'', ''APPLEGARTH'', ''APPLETON'', ''APRIL'', ''ARBURTHNOT'', ''ARBUTHNOT'', ''ARCHAMBAULT'', ''ARCHBOLD'', ''ARCHDEKEN/ARCHDEKIN'', ''ARCHER'', ''ARCHIBALD'', ''ARCONAITE'', ''ARDAGH'', ''ARDELL'', ''ARDIEL'', ''ARDIELL'', ''ARDILL'', ''ARDLE'', ''ARGUE'', ''ARHART'', ''ARIVISON'', ''ARKELL'', ''ARKENBRACK'', ''ARKLAND'', ''ARLESS'', ''ARLINGTON'', ''ARLSHOW'', ''ARMISTEAD'', ''ARMITAGE'', ''ARMON'', ''ARMOUR'', ''ARMS'', ''ARMSON'', ''ARMSTRONG'', ''ARMSTRONG/DOPP'', ''ARNER'', ''ARNISON'', ''ARNIT'', ''ARNOLD'', ''ARNOLD/CLOWES'', ''ARNOLDI'', ''ARNOT'', ''ARNOTE'', ''ARNOTT'', ''ARQUART'', ''ARROWSMITH'', ''ARTER/ARTIS'', ''ARTHUR'', ''ARTHUR/GIBSON'', ''ARTHURS'', ''ARVISON'', ''ASBURN'', ''ASH'', ''ASHBAUGH'', ''ASHBLAUGH'', ''ASHBOUGH'', ''ASHBRIDGE'', ''ASHBURN'', ''ASHBURRY'', ''ASHBY'', ''ASHFIELD'', ''ASHFORD'', ''ASHFORD/CROSS'', ''ASHLEY'', ''ASHMAN'', ''ASHTON'', ''ASHWELL'', ''ASHWORTH'', ''ASKEW'', ''ASKIN'', ''ASKINS'', ''ASKWITH'', ''ASLEY'', ''ASOPH'', ''ASPREY'', ''ASSEKINACK'', ''ASSELIN'', ''ASSELSBYNE/CORBAN'', ''ASSELSTINE'', ''ASSELSTINE/DAVY'', ''ASSELTINE'', ''ASTON'', ''ATCHESON'', ''ATCHISON'', ''ATKIN'', ''ATKINS'', ''ATKINSON'', ''ATTHILL'', ''ATTWELL'', ''ATTWOOD'', ''ATWOOD'', ''AUBERRY'', ''AUBERT'', ''AUGER'', ''AUGER/ANGER'', ''AUGUSTIN'', ''AUGUSTINE'', ''AULDJO'', ''AULDJO & MAITLAND'', ''AULDJO/McDONELL'', ''AULIN'', ''AULT'', ''AULT/LOUCKS'', ''AUMOND'', ''AUREY'', ''AURT'', ''AUSEM'', ''AUSMAN'', ''AUSTEN'', ''AUSTIN'', ''AUTIN'', ''AUTINO'', ''AVERILL'', ''AVERY'', ''AVISON'', ''AVRILL'', ''AWREY'', ''AXFORD'', ''AXHAM'', ''AYER'', ''AYLBOYNE'', ''AYLEN'', ''AYLMER'', ''AYLSWORTH'', ''AYLUM'', ''AYLWARD'', ''AYNGE'', ''AYOTT'', ''AYRES'', ''Abbot'', ''Abbott'', ''Abell'', ''Abrams/Lockwood'', ''Ackler/Warren'', ''Ackley/Meritt'', ''Acque'', ''Acre/Oswald'', ''Acre/Trowbridge'', ''Acres'', ''Acton'', ''Adair/Springsteen'', ''Adair/Ward'', ''Adam'', ''Adam/Mattice'', ''Adams'', ''Adams/McNeil'', ''Adams/Nudle'', ''Adams/Sage'', ''Adams/Snider'', ''Adams/Willson'', ''Adamson'', ''Affairs of the Canada Company with Government.'', ''African Baptist Church'', ''Agency of the Seigniory of Lauzon'', ''Agricultural Society'', ''Agricultural Society Middlesex'', ''Ahern'', ''Aidel'', ''Aird'', ''Airhash/Wager'', ''Airth'', ''Aitkin'', ''Aker/Jones'', ''Albrant/Statard'', ''Albrant/Stodard'', ''Albright'', ''Alderson'', ''Aldridge'', ''Algeor'', ''Alger'', ''Alger/Mosher'', ''Alger/Van Blarcom'', ''Algie'', ''Algire'', ''Algire/Van Blarcom'', ''Alguira/Jacocks'', ''Alguire/Loucks'', ''Alguire/Stallmayer'', ''Alguire/Wallisoe'', ''Allair'', ''Allan'', ''Allan/Wintermute'', ''Allbright/Albright'', ''Allen'', ''Allen/Mirick'', ''Allen/Vanalstine'', ''Allen/Watson'', ''Allison'', ''Allison/McTaggart'', ''Alphin'', ''Alston'', ''Altemas'', ''Althouse/Van Buren'', ''Alvord'', ''Alway'', ''Alyea'', ''Amalicite Indians'', ''Amey/Snyder'', ''Amory'', ''Anable'', ''Anderson'', ''Anderson/LaCompte'', ''Anderson/Lundy'', ''Anderson/Macdonell'', ''Anderson/Main'', ''Anderson/McCuen'', ''Anderson/Reid'', ''Anderson/Ruderback'', ''Anderson/Secord'', ''Andrew'', ''Andrews'', ''Andrews/Morrill'', ''Andrews/Mott'', ''Andrews/Wright'', ''Angell'', ''Anger/House'', ''Anger/McClintock'', ''Anglican Church and Church Wardens'', ''Anglican Presbyterians'', ''Angst'', ''Anguish/Lane'', ''Angus'', ''Annabel'', ''Annable'', ''Annable/Rowley'', ''Anneser/Gibson'', ''Annible/Wagner'', ''Ansley/Merrill'', ''Anson'', ''Anton/Molyneux'', ''Apostolic Missionary Church Wardens'', ''Applegarth'', ''Applicants for Clergy Reserves'', ''Applicants for free grants'', ''Applicants of Niagara'', ''Application for Deeds, township of Arthur.'', ''Application for lease of occupation of a certain block'', ''Application for the Township of Grantham'', ''Application for the concession of "Dewar Island"'', ''Application of settlers for survey of Somerset, Stanfold and Maddington'', ''Application to purchase Clergy Reserves in the County of Carleton.'', ''Application to purchase lands in the township of Inverness'', ''Application to purchase lot'', ''Application to purchase lots in Orillia'', ''Applications For the purchase of certain Clergy Reserved in Amaranth and Melancthon.'', ''Applications for Lot 31 in 13 Concession in Zorra Township'', ''Applications for grants'', ''Applications for land'', ''Applications for land script of officers'', ''Applications for mining'', ''Applications for purchase of Crown lands in Zone'', ''Applications of the municipality of St. Vincent'', ''Applications to buy lots'', ''Applications to purchase land in -'', ''Applications to purchase land in Sheffield'', ''Applications to purchase lands'', ''Applications to purchase lands in Gore of Somerset'', ''Applications to purchase lots 26 & 27'', ''Appointment of Inspectors of Land in Upper and Lower Canada'', ''Archar'', ''Archbishop of Quebec'', ''Archer'', ''Archibald'', ''Archiepiscopal Corporation of Quebec'', ''Ardill'', ''Argue'', ''Arkell'', ''Arkland'', ''Arkley'', ''Armatage'', ''Armour'', ''Armstrong'', ''Armstrong/Marlett'', "Armstrong/O''Neill", ''Armstrong/Welch'', ''Arner/Whittle'', ''Arnold'', ''Arnold/Munson'', ''Arnot'', ''Arrangements proposed by the Receiver General (Peter Russell) For sending of patents to the different districts.'', ''Arundel in relation to survey'', ''Ash'', ''Ashford/Mallery'', ''Ashley'', ''Ashman'', ''Ashton'', ''Askin'', ''Askin/Richardson'', ''Askins/Weatherhead'', ''Aspenwall'', ''Asselstine'', ''Asselstine/Lessard'', ''Asselstine/Watson'', ''Association of Bellechasse.'', ''Association of Bellechasse; relative to lands in Buckland.'', ''Asten/Snitzinger'', ''Atchison'', ''Atkins'', ''Atkinson'', ''Attorney General issue of a fiat dispensed with in certain cases.'', ''Attorney General re Church Land'', ''Attorney General respecting Survey near Niagara Falls.'', ''Attorney General respecting leases of the Crown and Clergy Reserves.'', ''Attorney General respecting surrender of certain patents in the township of Hope.'', "Attorney General''s report", ''Attorney General, opinion land set apart for Government House had been promised as an endowment for hospital.'', ''Attorney General, opinion respecting delivery of Patent Deeds.'', ''Auction'', ''Auger/Near'', ''Ault/Redmond'', ''Ault/Snider'', ''Ault/Thomb'', ''Ault/Weart'', ''Austin/Marr'', ''Authority for verification of outlines of Shawenagan'', ''Authority to advertise lots'', ''Authority to perform Surveys'', ''Avery'', ''Awrey'', ''Aynge'', ''Ayotte'', ''Ayrs'', ''Ayton'', ''Azereas'', ''BAAR'', ''BABBIT'', ''BABCOCK'', ''BABCOCK/AYNHART'', ''BABCOCK/CLEVELAND'', ''BABCOCK/COLE'', ''BABCOCK/COOK'', ''BABCOCK/HARPELL'', ''BABCOCK/SNIDER'', ''BABCUCK'', ''BABESK/COLE'', ''BABINEAU'', ''BABY'', ''BACH'', ''BACHAMAN'', ''BACHE'', ''BACHUS'', ''BACKEN/BACON'', ''BACKER'', ''BACKHOUSE'', ''BACKHOUSE/McMICHAEL'', ''BACKSTER'', ''BACKUS'', ''BACON'', ''BADDELEY'', ''BADDER'', ''BADGELY'', ''BADGEN'', ''BADGER'', ''BADGERON'', ''BADGLEY'', ''BADICHON'', ''BAGE'', ''BAGERO'', ''BAGG'', ''BAGLEY'', ''BAGSHAW'', ''BAGSLEY'', ''BAGWELL'', ''BAID'', ''BAIL/PRETEUS'', ''BAILEY'', ''BAILEY/BALEY'', ''BAILEY/DOVER'', ''BAILEY/WINTERMULER'', ''BAILIE'', ''BAILLIE'', ''BAILLY'', ''BAILY'', ''BAIN'', ''BAINES'', ''BAIRD'', ''BAIRDEN'', ''BAITS'', ''BAITSON'', ''BAKEN'', ''BAKEN/SMITH'', ''BAKER'', ''BAKER/CARPENTER'', ''BAKER/COOK'', ''BAKER/DAVY'', ''BAKER/FRANKLIN'', ''BAKER/FRANKS'', ''BAKER/MILLARD'', ''BAKER/McARTHUR'', ''BAKER/NOODLE'', ''BAKEWELL'', ''BALANCE'', ''BALCH'', ''BALDERSON'', ''BALDWIN'', ''BALDWIN/SERVICE'', ''BALDWIN/SHAW'', ''BALDWINS'', ''BALDWYN'', ''BALEY'', ''BALFOUR'', ''BALIE'', ''BALKE'', ''BALKWELL'', ''BALKWILL'', ''BALL'', ''BALL/CREEN'', ''BALL/CRYSLER'', ''BALLANTYNE'', ''BALLARD'', ''BALLENTINE'', ''BALM/DOCHSTADER'', ''BALMER'', ''BALTEEL'', ''BALTY'', ''BAMBER'', ''BAMBRIDGE/BAMBIDGE'', ''BAMFORD'', ''BANCK'', ''BANCRAFT'', ''BANCROFT'', ''BANCROFT/BENCROFT'', ''BANDER'', ''BANE'', ''BANER'', ''BANFIELD'', ''BANGHART'', ''BANGHART/BAUGHART'', ''BANGS'', ''BANKER'', ''BANKS'', ''BANNEAU'', ''BANNERMAN'', ''BANNERMAN/BOWERMAN'', ''BANNISTER'', ''BANNON'', ''BANTA'', ''BAPTIST'', ''BARABE'', ''BARBAR'', ''BARBARE'', ''BARBAREE'', ''BARBEAU'', ''BARBER'', ''BARBER/PELLET'', ''BARBERE'', ''BARBO'', ''BARBOUR'', ''BARBU'', ''BARCHKLY'', ''BARCHLE'', ''BARCKLEY'', ''BARCKLY'', ''BARCLAY'', ''BARCLEY'', ''BARCLEY/BARKELLEY'', ''BARDAN/CLARK'', ''BARDBURN'', ''BARDEN'', ''BARDIN/McNIEL'', ''BARDON'', ''BARDSLEE'', ''BAREFOOT'', ''BARGAR'', ''BARGER'', ''BARGER/FOSTER'', ''BARHLEY'', ''BARIL'', ''BARJARO'', ''BARK'', ''BARKALEY'', ''BARKELLY'', ''BARKER'', ''BARKHOUSE'', ''BARKLAY'', ''BARKLEY'', ''BARKLEY/AULT'', ''BARKLEY/GARLOUGH'', ''BARKLEY/MARKLEY'', ''BARKLEY/STRADER'', ''BARKLEY/UTTMAN'', ''BARLEY'', ''BARLINGER'', ''BARLOU/BAILIE'', ''BARLOW'', ''BARLOWE'', ''BARNAM'', ''BARNARD'', ''BARNELL'', ''BARNES'', ''BARNET'', ''BARNET/LANDON'', ''BARNETT'', ''BARNEY'', ''BARNHAM'', ''BARNHAM/DACHSTEDER'', ''BARNHARD'', ''BARNHART'', ''BARNHART/MARSH'', ''BARNHART/PRENTICE'', ''BARNHART/RAMBOUGH'', ''BARNHART/SILLS'', ''BARNHART/WAGER'', ''BARNIER'', ''BARNS'', ''BARNSIDE'', ''BARNUM'', ''BARNUM/CLARYSDALE'', ''BARNUM/HAME'', ''BARO'', ''BARON'', ''BARON/LAFRENIERE'', ''BARR'', ''BARRACLOUGH'', ''BARRAGER/LEE'', ''BARREL'', ''BARRELL'', ''BARRET'', ''BARRETT'', ''BARRETT/WALKER'', ''BARRETTE'', ''BARRICHER/STRADER'', ''BARRIE'', ''BARRIER'', ''BARRIGER'', ''BARRIL'', ''BARRITT'', ''BARRON'', ''BARROWMAN'', ''BARRY'', ''BART'', ''BARTCH'', ''BARTELEY'', ''BARTELL/CARSCALLEN'', ''BARTELS'', ''BARTEN'', ''BARTEN/DIXON'', ''BARTH'', ''BARTHE'', ''BARTHOL'', ''BARTHOLD'', ''BARTHOLEMEW'', ''BARTHOLOMAY'', ''BARTHOLOMEW'', ''BARTHOLOMEY'', ''BARTHY'', ''BARTLEMAN'', ''BARTLES'', ''BARTLES/McDOUGAL'', ''BARTLET'', ''BARTLET, SMITH & others'', ''BARTLETT'', ''BARTLETT/DINGMAN'', ''BARTLEY'', ''BARTLY'', ''BARTLY/NICHOLSON'', ''BARTON'', ''BARTON/BARTIN'', ''BARTON/BOURDETT'', ''BARTON/McCORD'', ''BARTON/SMADES'', ''BARTOW'', ''BARTRON'', ''BARWICK'', ''BARWIS'', ''BARY'', ''BASKERRIL/BASKERVIL'', ''BASKERVILLE'', ''BASS'', ''BASS/AMEY'', ''BASSELL'', ''BASSET'', ''BASSETT'', ''BASSETT/HOPPLE'', ''BASSEY'', ''BASSEY/SMITH'', ''BASTABLE'', ''BASTEDO'', ''BASTIEN'', ''BATE'', ''BATEMAN'', ''BATER'', ''BATES'', ''BATES/MATTICE'', ''BATESON'', ''BATHGATE'', ''BATHIE'', ''BATHURST'', ''BATIE'', ''BATNAN'', ''BATTAM'', ''BATTAY'', ''BATTER'', ''BATTER/BEDFORD'', ''BATTER/POTTER'', ''BATTER/WRIGHT'', ''BATTERHAM'', ''BATTERS'', ''BATTERSBY'', ''BATTGER'', ''BATTIE'', ''BATTLE'', ''BATTY'', ''BATTY/LOYST'', ''BATY'', ''BAUBIN'', ''BAUGH'', ''BAUGHANAN'', ''BAUMES'', ''BAUMSTARK'', ''BAVENS'', ''BAXTER'', ''BAXTER/BINDER'', ''BAXTER/MILLER'', ''BAXTER/OATMAN'', ''BAXTER/SHIPMAN'', ''BAXTON'', ''BAYCROFT'', ''BAYER'', ''BAYEUR'', ''BAYLEY'', ''BAYLY'', ''BAYMAN'', ''BAYMON'', ''BAYNE'', ''BAYNES'', ''BEACH'', ''BEACH/FREO'', ''BEACH/SEELYE'', ''BEACHER'', ''BEACON'', ''BEADFORD'', ''BEADLE'', ''BEADSTEAD'', ''BEAGLE'', ''BEAKON'', ''BEAL'', ''BEALE'', ''BEALEY'', ''BEALL'', ''BEALSEY'', ''BEAM'', ''BEAM/MAY'', ''BEAM/MILLS'', ''BEAMAN'', ''BEAMER'', ''BEAMER/DOYLE'', ''BEAMISH'', ''BEAN'', ''BEANE'', ''BEARD'', ''BEARD/STRADER'', ''BEARDSLEY'', ''BEARDSLEY/BANKER'', ''BEARE'', ''BEARMAN'', ''BEARSS'', ''BEARSS/STEEL'', ''BEASLEY'', ''BEASTED/FRYMIRE'', ''BEATIE'', ''BEATON'', ''BEATRIE'', ''BEATSKE'', ''BEATSY'', ''BEATTEY'', ''BEATTIE'', ''BEATTIE/SCRAM'', ''BEATTON'', ''BEATTY'', ''BEATY'', ''BEATYS'', ''BEAUBIEN'', ''BEAUCHAMP'', ''BEAUDET'', ''BEAUDRY'', ''BEAUGRAND'', ''BEAUMON'', ''BEAUPORT'', ''BEAUPR?'', ''BEAUPR?/BOUPRIE'', ''BEAUPRE'', ''BEAUREGARD'', ''BEAVIS'', ''BEBEE'', ''BEBEE/GROOMS'', ''BEBEE/SEELY'', ''BECHER'', ''BECHSTED'', ''BECHTEL'', ''BECHTELL'', ''BECK'', ''BECKER'', ''BECKER/BAKER'', ''BECKERSTAFF'', ''BECKERTON'', ''BECKETT'', ''BECKETT/BOWMAN'', ''BECKFORD'', ''BECKON'', ''BECKON/BACON'', ''BECKON/SKINNER'', ''BECKOR'', ''BECKWITH'', ''BECKWITH/SILLS'', ''BECKWORTH'', ''BEDAL'', ''BEDARD'', ''BEDELL'', ''BEDFORD'', ''BEDICK/BEDSTED'', ''BEDSTEAD/CASSELMAN'', ''BEDSTED'', ''BEDSTED/HANESTER'', ''BEE'', ''BEEBE'', ''BEEBE/HART'', ''BEECH'', ''BEECHER'', ''BEEDELL'', ''BEEDLE'', ''BEEMAN'', ''BEEMER'', ''BEER'', ''BEERS'', ''BEESON'', ''BEETEN'', ''BEETHAM'', ''BEETON'', ''BEEVAN'', ''BEGG'', ''BEGG/STAG'', ''BEGHAN'', ''BEGUE'', ''BEIKIE'', ''BEITH'', ''BELANGER'', ''BELCHER'', ''BELCOUR'', ''BELFORD'', ''BELISLE'', ''BELKNAP'', ''BELL'', ''BELL/CARSCALLEN'', ''BELL/FINKLE'', ''BELL/FRALIC'', ''BELL/HARNES'', ''BELL/HITCHCOCK'', ''BELL/MERKLEY'', ''BELL/PALMER'', ''BELLAIRE'', ''BELLAMY'', ''BELLAN'', ''BELLANT'', ''BELLARD'', ''BELLECK'', ''BELLECOUR'', ''BELLEFEUILLE'', ''BELLEPERCHE'', ''BELLER'', ''BELLERE'', ''BELLINGER'', ''BELLINGER/OBERHOLD'', ''BELLINGER/SLINGERLAND'', ''BELLINGHAM'', ''BELLIVEAU'', ''BELLOWS'', ''BELLUE'', ''BELNAP'', ''BELNAP/MASTERS'', ''BELONGE'', ''BELROSE'', ''BELTON'', ''BELTZ'', ''BELYEA'', ''BEM'', ''BEMAN'', ''BEMER'', ''BEN'', ''BENAC'', ''BENCH'', ''BENCROFT'', ''BENDER'', ''BENDER/COLLINGER'', ''BENDER/WAIT'', ''BENDER/WOOD'', ''BENECKER'', ''BENEDICK'', ''BENEDICT'', ''BENEDICT/COLE'', ''BENEDICT/COVERLY'', ''BENEDICT/CRAWFORD'', ''BENEDICT/DAFOE'', ''BENEKER'', ''BENERDICT'', ''BENERE'', ''BENETEAU'', ''BENETEAU/LABALEINE'', ''BENFORD'', ''BENHAM'', ''BENINGER'', ''BENITEAU'', ''BENJAMIN'', ''BENKER'', ''BENKLEY'', ''BENLEY'', ''BENN'', ''BENN/BEHN'', ''BENN/FITCHETE'', ''BENNER'', ''BENNER/HOUCK'', ''BENNET'', ''BENNETT'', ''BENNIER'', ''BENNINGER'', ''BENNINGHAM'', ''BENNINGTON'', ''BENNITT'', ''BENNS'', ''BENO'', ''BENOIT'', ''BENSACK'', ''BENSEN'', ''BENSON'', ''BENSON/SAMSON'', ''BENSON/SCHARAMAHAN'', ''BENSON/SHEA'', ''BENTLEY'', ''BENTON'', ''BENVILL'', ''BENVILLE'', ''BERCKLY'', ''BERCSY'', ''BERCY'', ''BERCZY'', ''BERDAN'', ''BERDAN/McKENZIE'', ''BERDEBT/STEEL'', ''BERDEN'', ''BERFORD'', ''BERGAR'', ''BERGEN'', ''BERGER'', ''BERGER/TINKLE'', ''BERGIN'', ''BERGON'', ''BERHAM'', ''BERINGER'', ''BERISFORD'', ''BERKEY'', ''BERKLES'', ''BERN'', ''BERN/BEN'', ''BERNARD'', ''BERNARD/LARIVI?RE'', ''BERNARDIN'', ''BERNASH'', ''BERNER'', ''BERNIER'', ''BERNINGER'', ''BERO'', ''BERREL'', ''BERRIE'', ''BERRY'', ''BERTHEAUME'', ''BERTHELET'', ''BERTHIAUME'', ''BERTLE'', ''BERTLES'', ''BERTRAND'', ''BERTRAND/BERTRAN'', ''BERTRON'', ''BERUBE'', ''BERY'', ''BESENERE'', ''BESENET'', ''BESSE'', ''BESSEE'', ''BESSEY'', ''BESSEY/BESSY'', ''BESSIGER/BOSSIGER'', ''BEST'', ''BESWETHERICK'', ''BESWICK'', ''BESWICKE'', ''BETHAN'', ''BETHUME'', ''BETHUNE'', ''BETHUNE/SMITH'', ''BETRON'', ''BETRY'', ''BETTERIDGE'', ''BETTES'', ''BETTIE'', ''BETTIS'', ''BETTON'', ''BETTRIDGE'', ''BETTS'', ''BETTUNE'', ''BETTY'', ''BETTY/ORSER'', ''BETTYS'', ''BETZNER'', ''BEVANS'', ''BEVEAY'', ''BEVEN'', ''BEVERIDGE'', ''BEVIS'', ''BEWELL'', ''BEWLEY'', ''BEYERS'', ''BEYNON'', ''BEZER'', ''BICE'', ''BICE/UMPHRY'', ''BICKELL'', ''BICKERSON'', ''BICKET'', ''BICKFORD'', ''BICKHAM'', ''BIDOUT'', ''BIDWELL'', ''BIEKIE'', ''BIGCRAFT'', ''BIGCROFT'', ''BIGELOW'', ''BIGELOW/COLE'', ''BIGFORD'', ''BIGGAR'', ''BIGGAR/COLTMAN'', ''BIGGER'', ''BIGGINS'', ''BIGHAM'', ''BIGNELL'', ''BIGRAS'', ''BILES'', ''BILLEDEAUX'', ''BILLIDOE'', ''BILLIET'', ''BILLING'', ''BILLINGS'', ''BILLINGS/DAVY'', ''BILLINGS/GRAHAM'', ''BILLINGSBY'', ''BILLOPP'', ''BILLS'', ''BILOW/ALGIRE'', ''BILROSE/BONESTEEL'', ''BILTON'', ''BINAULT'', ''BINGHAM'', ''BINGHAM/BOWMAN'', ''BINGHAM/HUFFMAN'', ''BINGHAM/LAPTHORNE'', ''BINGLE'', ''BINIKER/BINICHER'', ''BINKLY'', ''BIRCH'', ''BIRD'', ''BIRD/CARSON'', ''BIRDS'', ''BIRDSALL'', ''BIRDSELL'', ''BIRMINGHAM'', ''BIRNEY'', ''BIRNIE'', ''BIRRELL'', ''BISHOP'', ''BISLEY/BISBEE'', ''BISSEL'', ''BISSELL'', ''BISSELL/CLARK'', ''BISSET'', ''BISSETT'', ''BISSIL'', ''BISSONET'', ''BISSONNETTE'', ''BLACK'', ''BLACK/CROSS'', ''BLACK/PICKLE'', ''BLACKBURN'', ''BLACKBURN/BAKER'', ''BLACKBURNE'', ''BLACKBURNS'', ''BLACKELY'', ''BLACKER'', ''BLACKER/MEYERS'', ''BLACKIE'', ''BLACKLEY'', ''BLACKLOCK'', ''BLACKLY'', ''BLACKLY/COLE'', ''BLACKMAN'', ''BLACKMOR'', ''BLACKNEY'', ''BLACKSOCK'', ''BLACKSTOCK'', ''BLACKWAY'', ''BLACKWELL'', ''BLACKWOOD'', ''BLADGEN'', ''BLADGET'', ''BLAIR'', ''BLAIR/ORO'', ''BLAIR/YOUNG'', ''BLAKE'', ''BLAKE/ADAMS'', ''BLAKELEY'', ''BLAKELY'', ''BLAKEMORE'', ''BLAKENEY'', ''BLAKER'', ''BLAKER/BLEECKER'', ''BLAKER/MEYERS'', ''BLAKESBY'', ''BLAKESLEY'', ''BLAKEY'', ''BLAKLEY'', ''BLAKLEY/BENNIT'', ''BLAKLY'', ''BLAKNEY'', ''BLAM'', ''BLANCHARD'', ''BLANCHARD/RUSH'', ''BLANCHER'', ''BLANCHET'', ''BLANCHFILL'', ''BLANCO'', ''BLAND/CAMPBELL'', ''BLANEY'', ''BLANN'', ''BLASDELL'', ''BLAUVELT'', ''BLAYLOCK'', ''BLEACHUM'', ''BLEAKLY'', ''BLEAM'', ''BLEECKER'', ''BLEEKER'', ''BLEEKER/MEYERS'', ''BLENEY'', ''BLETCHER'', ''BLEVINS'', ''BLEW'', ''BLEWETT'', ''BLEWITT'', ''BLEZARD'', ''BLINN'', ''BLISS'', ''BLIZARD'', ''BLODGET/WHEELER'', ''BLODGETT'', ''BLOM'', ''BLOOM'', ''BLOOMFIELD'', ''BLOOR'', ''BLOSDOLL'', ''BLOUNT'', ''BLOVIN'', ''BLUCHERS'', ''BLUCKER'', ''BLUE'', ''BLUEMAN'', ''BLUETT'', ''BLUM'', ''BLUNDEN'', ''BLUNT'', ''BLYTH'', ''BLYTH/CORBMAN'', ''BOALAN'', ''BOARDMAN/WINTERS'', ''BOBIER'', ''BOCAUT'', ''BOCHET'', ''BOCHUS'', ''BOCKUS'', ''BOCKUS/SHAVER'', ''BODDINGTON'', ''BODDY'', ''BODIN'', ''BODINE'', ''BODY'', ''BOG'', ''BOGARD'', ''BOGART'', ''BOGART/BENEDICT'', ''BOGERT'', ''BOGGS'', ''BOGUE'', ''BOHAN'', ''BOHANNON'', ''BOICE'', ''BOICE/AMEY'', ''BOICE/CUSHMAN'', ''BOICE/SCHRIVER'', ''BOID'', ''BOILS'', ''BOIRE'', ''BOISE'', ''BOISEAU'', ''BOISMIER'', ''BOISNISER'', ''BOISVERT'', ''BOITON'', ''BOLAND'', ''BOLBY'', ''BOLDUC'', ''BOLE'', ''BOLES'', ''BOLLARD'', ''BOLSTER'', ''BOLTE'', ''BOLTON'', ''BOLTON/BOULTON'', ''BOLTON/BUELL'', ''BOLWELL'', ''BOMBARDIER'', ''BONBERGIER'', ''BOND'', ''BOND/PATTERSON'', ''BONDIE'', ''BONDRAY'', ''BONDY'', ''BONE'', ''BONEAU'', ''BONER'', ''BONESTEAL'', ''BONESTEEL'', ''BONESTIEL'', ''BONGARD'', ''BONKER'', ''BONKER/CORTEN'', ''BONNELL'', ''BONNER'', ''BONNET'', ''BONNETT'', ''BONNEY'', ''BONNYCASTLE'', ''BONSACK'', ''BONSER'', ''BONSTED'', ''BONSTEEL'', ''BONSTEEL/CONWAY'', ''BONTER'', ''BONTER/DEMPSEY'', ''BONTER/PEAK'', ''BONYARD'', ''BOODLE'', ''BOODY'', ''BOOK'', ''BOOMER'', ''BOORDEN'', ''BOORMAN'', ''BOOS'', ''BOOTH'', ''BOOTH/DALY'', ''BOOTH/NORMAN'', ''BOOTHROYD'', ''BORCKHOLDER'', ''BORDEN'', ''BORDEN/SEELEY'', ''BORGHONDER'', ''BORING'', ''BORKHOLDER'', ''BORLAND'', ''BORNUM'', ''BORPY'', ''BORRELLY'', ''BORRIE'', ''BORUGHTON/DENNIS'', ''BOSQUET'', ''BOSQUET/COX'', ''BOSS'', ''BOSS/AINSLEY'', ''BOSSE'', ''BOSSENIE'', ''BOSSET/SMITH'', ''BOSSU/LEONNAIS'', ''BOSTICK'', ''BOSTLY'', ''BOSTON'', ''BOSTWICK'', ''BOSTWICKS'', ''BOSWELL'', ''BOSWICK'', ''BOTHWELL'', ''BOTSFORD'', ''BOTTERILL'', ''BOTTOM'', ''BOTTON'', ''BOTTON/CURRY'', ''BOTTUM'', ''BOUCH'', ''BOUCH?'', ''BOUCHE'', ''BOUCHER'', ''BOUCHER DE BOUCHERVILLE'', ''BOUCHERVILLE'', ''BOUCHETTE'', ''BOUCHETTE/BERTHELET'', ''BOUCHNER'', ''BOUCK'', ''BOUCK/DILLABAUGH'', ''BOUCKE'', ''BOUDETTE'', ''BOUDREAU'', ''BOUET'', ''BOUFFARD'', ''BOUGENER'', ''BOUGHMAN'', ''BOUGHNER'', ''BOUGHNER/GLOVER'', ''BOUGHNER/MATHEWS'', ''BOUGHSLAUCH'', ''BOUKE'', ''BOUKER'', ''BOULANGER'', ''BOULBY'', ''BOULD'', ''BOULGER'', ''BOULLARD'', ''BOULT'', ''BOULTER/PEACK'', ''BOULTON'', ''BOULTON/BELLANY'', ''BOULTON/BOTTON'', ''BOULTON/CURRY'', ''BOULTON/ELLIOT'', ''BOULTON/ELLIOTT'', ''BOULTON/JONES'', ''BOULTON/REDMAN'', ''BOURASSA'', ''BOURBONNIER'', ''BOURCHIER'', ''BOURDETT'', ''BOURDETT/STEEL'', ''BOURDETTE'', ''BOURGET'', ''BOURJOTTE/LESPERANCE'', ''BOURK'', ''BOURKE'', ''BOURNS'', ''BOUS'', ''BOUSLAUGH'', ''BOUSSEY'', ''BOUTHELLIER'', ''BOUTHILLIER'', ''BOUTILLET'', ''BOUTILLIER'', ''BOUTIN'', ''BOVARD'', ''BOWAN'', ''BOWBEER'', ''BOWBEER/BREWER'', ''BOWBIER'', ''BOWDEN'', ''BOWEN'', ''BOWEN/CARMAN'', ''BOWEN/CULBERTSON'', ''BOWEN/DEMOREST'', ''BOWEN/DIMONT'', ''BOWEN/KIMMERLY'', ''BOWEN/PALMER'', ''BOWEN/PORT'', ''BOWEN/RISELAY'', ''BOWEN/SHAVER'', ''BOWENS'', ''BOWER'', ''BOWERMAN'', ''BOWERMAN/MORDEN'', ''BOWERS'', ''BOWES'', ''BOWIE'', ''BOWIE/BUIES'', ''BOWKER'', ''BOWKET'', ''BOWKETT'', ''BOWLBY'', ''BOWLES'', ''BOWLEY'', ''BOWLING'', ''BOWMAN'', ''BOWMAN/BELLINGER'', ''BOWMAN/BOUKE'', ''BOWMAN/CLARK'', ''BOWMAN/DEFOREST'', ''BOWMAN/MATLACK'', ''BOWN'', ''BOWNEY'', ''BOWRON'', ''BOWSER/BOWSEN'', ''BOWSFIELD'', ''BOXALL'', ''BOYCE'', ''BOYCE/HARTMAN'', ''BOYD'', ''BOYD/McINTIRE'', ''BOYDE'', ''BOYEA'', ''BOYER'', ''BOYER/COMER'', ''BOYERS'', ''BOYES'', ''BOYL'', ''BOYLAN'', ''BOYLE'', ''BOYLES'', ''BOYS'', ''BOYSEE'', ''BOYTON'', ''BRACK'', ''BRACKBILL'', ''BRACKEN'', ''BRACKENREED'', ''BRACKENRIDGE'', ''BRACKETT'', ''BRACKIN'', ''BRACKLEY'', ''BRADBENT'', ''BRADBURN'', ''BRADBURNE'', ''BRADEY'', ''BRADFORD'', ''BRADFORD/SELEE'', ''BRADISH'', ''BRADLEY'', ''BRADLEY/MOOR'', ''BRADLY'', ''BRADNER'', ''BRADNOCK'', ''BRADSHAW'', ''BRADSHAW/COLLINS'', ''BRADSHAW/DIAMOND'', ''BRADSHAW/LARUE'', ''BRADSHAW/VANDERLISS'', ''BRADSTEAD/AULT'', ''BRADT'', ''BRADT/BENDER'', ''BRADT/HAINER'', ''BRADT/YOUNG'', ''BRADWELL'', ''BRADY'', ''BRAGG'', ''BRAILY'', ''BRAKENRIDGE'', ''BRALY'', ''BRAMORSKEY'', ''BRAMORSKY'', ''BRAMOUSKY'', ''BRANAGAN'', ''BRANAN'', ''BRANDON'', ''BRANDT'', ''BRANDTS'', ''BRANDTS/BRANTS'', ''BRANDYMAN'', ''BRANDYMANN'', ''BRANIGAN'', ''BRANNAGH'', ''BRANNAH'', ''BRANNAN'', ''BRANNEGAN'', ''BRANNEN'', ''BRANNON'', ''BRANT'', ''BRANT/DESERONTWE'', ''BRANTFORD/REED'', ''BRANTS'', ''BRASH'', ''BRASS'', ''BRASS/MATTICE'', ''BRASS/MATTIN'', ''BRASS/ROUSEHORNE'', ''BRATH'', ''BRATT'', ''BRAUNCIS'', ''BRAY'', ''BRAY/ANDERSON'', ''BRAYLEY'', ''BRAYLEY/LAW'', ''BRAYTON'', ''BRAYTON/WEART'', ''BRAZEY'', ''BRAZIL'', ''BRAZLE'', ''BREAD'', ''BREADEN'', ''BREADNER'', ''BREADON'', ''BREADY'', ''BREAK'', ''BREAKENRIDGE'', ''BREAKENRIGHT'', ''BREAKY'', ''BREARLY'', ''BREATHWAITE'', ''BREAULT'', ''BREBNER'', ''BRECKAN'', ''BRECKON'', ''BREEN'', ''BREGARE'', ''BREINICH'', ''BREMBLE'', ''BREMICH'', ''BREMMER'', ''BREMNER'', ''BRENAGH'', ''BRENAN'', ''BRENEMAN'', ''BRENNAN'', ''BRENNAN/PETERSON'', ''BRENNEN'', ''BRENT'', ''BRENTON'', ''BRESEE'', ''BRESNAHAN'', ''BRESSE'', ''BRESSEE'', ''BRETHWAITH'', ''BRETNER'', ''BRETT'', ''BREWER'', ''BREWERTON'', ''BREZE'', ''BREZEE'', ''BRIAN'', ''BRICE'', ''BRICK'', ''BRICKER'', ''BRICKMAN'', ''BRICOT'', ''BRIDGAM'', ''BRIDGE'', ''BRIDGEFORD'', ''BRIDGEMAN'', ''BRIDGES'', ''BRIDGEWATER'', ''BRIDGLAND'', ''BRIDGMAN'', ''BRIEN'', ''BRIENLY'', ''BRIERLY'', ''BRIGGS'', ''BRIGHAM'', ''BRIGHAM/BREWSTER'', ''BRIGHT'', ''BRIGS'', ''BRIKIE'', ''BRILL'', ''BRILLARD'', ''BRILLINGER'', ''BRIMMER'', ''BRINDLEY'', ''BRINK'', ''BRINK/BRADT'', ''BRINKERHOOF'', ''BRINKMAN'', ''BRINLEY'', ''BRINSMEAD'', ''BRINTNALL'', ''BRIODY'', ''BRISARD'', ''BRISBOIS'', ''BRISCO'', ''BRISCOE'', ''BRISCOE/BELL'', ''BRISLANE'', ''BRISON'', ''BRISTOE'', ''BRISTOL'', ''BRITAIN'', ''BRITAIN/BRITTIN'', ''BRITTON'', ''BROAD'', ''BROADBANT'', ''BROADBENT'', ''BROADIE'', ''BROADWICK'', ''BROCK'', ''BRODAY'', ''BRODDY'', ''BRODEUR'', ''BRODIE'', ''BRODRICK'', ''BRODY'', ''BROEFFEL/BREEFFEL'', ''BROFFEL'', ''BROFFEY'', ''BROGDEN'', ''BROGDIN'', ''BROGDON'', ''BROMLY/ROSS'', ''BROOK'', ''BROOKE'', ''BROOKER'', ''BROOKES'', ''BROOKING'', ''BROOKS'', ''BROOKS/BARRETT'', ''BROOKS/BENEDICT'', ''BROOKS/COCKMAN'', ''BROOKS/COLE'', ''BROOKS/CORBMAN'', ''BROOKS/MEDDAUGH'', ''BROOKS/TURGERSON'', ''BROOKSBANK'', ''BROOME'', ''BROOMHEAD'', ''BROONER'', ''BROPHEY'', ''BROPHY'', ''BROSE'', ''BROTHERS'', ''BROTHERSTON'', ''BROUGH'', ''BROUGHMAN'', ''BROUGHTON'', ''BROUK'', ''BROUSE'', ''BROUSE/ADAMS'', ''BROUSE/BOUCK'', ''BROUSE/CARMAN'', ''BROUSE/KARMAN'', ''BROUSE/PARLOW'', ''BROUSE/SHAVER'', ''BROUSTER'', ''BROWCE'', ''BROWCE/COON'', ''BROWDERS'', ''BROWER'', ''BROWN'', ''BROWN/COON'', ''BROWN/FAIRFIELD'', ''BROWN/FERGUSON'', ''BROWN/FREEMAN'', ''BROWN/HAWEN'', ''BROWN/HURLBURT'', ''BROWN/JUDSON'', ''BROWN/MATICE'', ''BROWN/MCILMOYLE'', ''BROWN/MCLEAN'', ''BROWN/MOTT'', ''BROWN/MULLOY'', ''BROWN/McINTIRE'', ''BROWN/NICHOLSON'', ''BROWN/ROBBINS'', ''BROWN/SNIDER'', ''BROWNE'', ''BROWNELL'', ''BROWNELL/CLINE'', ''BROWNELL/COUNTRY'', ''BROWNELL/SERVICE'', ''BROWNELLE/BROWNHILL'', ''BROWNLEE'', ''BROWNLEY'', ''BROWNLIE'', ''BROWNLY'', ''BROWNSON'', ''BROWNSON/HAFFMAN'', ''BROWNWICK'', ''BROWSE'', ''BRUCE'', ''BRUCE/DAVY'', ''BRUCK'', ''BRUEN'', ''BRUGERE'', ''BRULE'', ''BRUMGET'', ''BRUMIGEM'', ''BRUMIGEM/BRUMIGEN'', ''BRUMWELL'', ''BRUND'', ''BRUNDAGE'', ''BRUNDAGE/YOUREX'', ''BRUNDEGE'', ''BRUNDICH'', ''BRUNDIDGE'', ''BRUNDIGE'', ''BRUNDREDGE'', ''BRUNDRIDGE'', ''BRUNELLE'', ''BRUNER'', ''BRUNGER'', ''BRUNHERD'', ''BRUNK'', ''BRUNSON'', ''BRUSH'', ''BRUSH/CURRY'', ''BRUSSEAUX'', ''BRUSTER'', ''BRUTON'', ''BRUTON/FINNERON'', ''BRYAN'', ''BRYANS'', ''BRYANT'', ''BRYANT/McLEAN'', ''BRYCE'', ''BRYDEN'', ''BRYDON'', ''BRYDONE'', ''BRYERS'', ''BRYNES'', ''BRYSON'', ''BUCANNAN/BUCHANNAN/BUCKANNAN'', ''BUCH'', ''BUCHAM'', ''BUCHAN'', ''BUCHANAN'', ''BUCHANAN & BARKER'', ''BUCHANAN/CRAWFORD'', ''BUCHANAN/CRYSLER'', ''BUCHANNAN'', ''BUCHLEY'', ''BUCHNER'', ''BUCHNER/BOUGENER'', ''BUCHNER/BOUGHNER'', ''BUCHNER/BOUGHNER/AINSLEY'', ''BUCHNER/BOUGNER'', ''BUCHNER/BUGHGONER'', ''BUCHNER/GRAHAM'', ''BUCK'', ''BUCK/BALDWIN'', ''BUCK/BOUCH'', ''BUCK/HAGERMAN'', ''BUCK/HAGGART'', ''BUCK/ORSER'', ''BUCK/SNUKE'', ''BUCKBOROUGH'', ''BUCKHAM'', ''BUCKHANAN/McDONELL'', ''BUCKLEY'', ''BUCKLEY/THOMAS'', ''BUCKNER'', ''BUCKNER/McGAN'', ''BUDD'', ''BUELL'', ''BUELL/BAUCKE'', ''BUELL/SHERWOOD'', ''BUFFA'', ''BUGBEE'', ''BUGG'', ''BUGGENER'', ''BUGGIE'', ''BUGIN'', ''BUGLE LE'', ''BUIE'', ''BUIS'', ''BUIS/AMEY'', ''BUIS/BICE'', ''BUISE'', ''BUKER/BAKER'', ''BUKER/FRENCH'', ''BULEY'', ''BULGER'', ''BULHAM'', ''BULKELEY'', ''BULL'', ''BULL/CLEMENT'', ''BULLARD'', ''BULLEN'', ''BULLER'', ''BULLES'', ''BULLESS'', ''BULLEY'', ''BULLIN'', ''BULLISS'', ''BULLIVET'', ''BULLIVIT'', ''BULLOCK'', ''BULLY'', ''BULMAN'', ''BULMER'', ''BULSON'', ''BULTEEL'', ''BUMP'', ''BUNBURRY'', ''BUNBURY'', ''BUNCE'', ''BUNDIGE'', ''BUNFIELD'', ''BUNKER'', ''BUNN'', ''BUNTING'', ''BUNTING/MARCKLE'', ''BUNTON'', ''BUNYAN'', ''BURBANK'', ''BURBANKS'', ''BURBY'', ''BURCH'', ''BURCH/BIRCH'', ''BURCHELL'', ''BURCHILL'', ''BURCK'', ''BURCKHOLDER'', ''BURD'', ''BURDEN'', ''BURDETT'', ''BURDETT/BOURDETT'', ''BURDICK'', ''BURDOCK'', ''BURDON'', ''BURDS'', ''BURDSELL'', ''BURDSLEE'', ''BUREAU'', ''BURGAR'', ''BURGART'', ''BURGER'', ''BURGES'', ''BURGESS'', ''BURGHAM'', ''BURGHER'', ''BURGHMANN'', ''BURGIS'', ''BURGISS'', ''BURGOINE'', ''BURGOYNE'', ''BURK'', ''BURKE'', ''BURKET'', ''BURKHOLDER'', ''BURLEIGH'', ''BURLEY'', ''BURLEY/BALDWIN'', ''BURLEY/LEE'', ''BURLEY/SNIDER'', ''BURLINGAME'', ''BURLINGHAM/TANKEN'', ''BURN'', ''BURNELL'', ''BURNES'', ''BURNET'', ''BURNET/VANARDEN'', ''BURNETT'', ''BURNETT/BAKER'', ''BURNETT/COLLINGWOOD'', ''BURNHAM'', ''BURNHAM/AUGER'', ''BURNHAM/CLAWSON'', ''BURNHAM/FERGUSON'', ''BURNHAM/HUFF'', ''BURNITT'', ''BURNS'', ''BURNS/BYRNE'', ''BURNSIDE'', ''BURNSIDE/BROUSE'', ''BURON'', ''BURR'', ''BURRELL'', ''BURRETT'', ''BURRILL'', ''BURRISON'', ''BURRISON/JOHNSTOWN'', ''BURRIT'', ''BURRIT/DULMAGE'', ''BURRITT'', ''BURROWES'', ''BURROWS'', ''BURRY'', ''BURST'', ''BURSTELL'', ''BURT'', ''BURTCH'', ''BURTCH/BRADT'', ''BURTCH/CORLIS'', ''BURTCH/SEHRAM'', ''BURTON'', ''BURTON/ALGUIRE'', ''BURWELL'', ''BURWELL/BENEDICT'', ''BURWELL/BISSELL'', ''BURY'', ''BUSCH'', ''BUSH'', ''BUSH/CROWDER'', ''BUSHBY'', ''BUSHEY'', ''BUSHNELL'', ''BUSNEY'', ''BUSSELL'', ''BUSVAK'', ''BUTCH'', ''BUTCHART'', ''BUTCHER'', ''BUTH'', ''BUTLER'', ''BUTLER/BERRIE'', ''BUTLER/BOURDETT'', ''BUTLER/COLTMAN'', ''BUTLER/CROOKS'', ''BUTLER/DAYTON'', ''BUTLER/MATTHEWS'', ''BUTLER/SLEMENT'', ''BUTLER/STEVENS'', ''BUTLER/WOOD'', ''BUTTER'', ''BUTTERFIELD'', ''BUTTERFIELD/CLARK'', ''BUTTERWORTH'', ''BUTTLER'', ''BUTTON'', ''BUTTS'', ''BUXTON'', ''BUZBY'', ''BY'', ''BYAN'', ''BYARD'', ''BYARS'', ''BYCRAFT'', ''BYER'', ''BYERLAY'', ''BYERS'', ''BYRNE'', ''BYRNES'', ''BYRNS'', ''Babcock'', ''Babcock/Lewis'', ''Babcock/Ransier'', ''Babcock/Scott'', ''Babcock/Shannon'', ''Babcock/Smith'', ''Babcock/Wannamaker'', ''Babcock/Wannomaker'', ''Babcock/Weese'', ''Babcook/Scott'', ''Bachman'', ''Bacon'', ''Badcock'', ''Badger'', ''Badgerow'', ''Badgley'', ''Badichon/Labadie'', ''Badwin'', ''Bag'', ''Bagley'', ''Baice/Langhorn'', ''Bailey'', ''Bailey/McCarthy'', ''Bailey/Sawyers'', ''Baily'', ''Bain'', ''Baine'', ''Baker'', ''Baker/Miller'', ''Baker/Skinner'', ''Baker/Smith'', ''Baker/Snyder'', ''Baker/Wade'', ''Baker/Waldorff'', ''Baker/Young'', ''Balance of land'', ''Balderson'', ''Baldwin'', ''Baldwin/Webster'', ''Bale'', ''Ball'', ''Ball/Long'', ''Ball/Servos'', ''Ball/Sweep'', ''Ball/Walker'', ''Ball/Wilkeson'', ''Balmain/McDonell'', ''Balton/Rose'', ''Bamber'', ''Banal Grist Mill'', ''Banal Mill'', ''Bannerman'', ''Bannister'', ''Banta/Vanwyck'', ''Banter'', ''Baptist Church'', ''Baptist Church at Amherstburg.'', ''Baptist Society'', ''Baptists of the township of Artemesia.'', ''Barber'', ''Barckley/Outerkirk'', ''Barckly/Markle'', ''Barclay'', ''Barclough'', ''Barden'', ''Barhigt'', ''Bark/Lyons'', ''Barkely/Shell'', ''Barkley/Whiteker'', ''Barkley/Witticker'', ''Barkly/Marcellia/Massales'', ''Barlow'', ''Barnard/Nichols'', ''Barnard/Randell'', ''Barnard/Thomas'', ''Barnard/Walker'', ''Barnet/Landon'', ''Barnet/Rood'', ''Barnhart/Gallingher'', ''Barnhart/MacDaniel'', ''Barnhart/Marks'', ''Barnhart/Rombough'', ''Barnhart/Sills'', ''Barnhart/Stinehoof'', ''Barnhart/Stringer'', ''Barnhart/Tilton'', ''Barnheart/Runion'', ''Barnheart/Waldenbarger'', ''Barnrom'', ''Barns'', ''Barr'', ''Barrie'', ''Barrie & Kempenfeldt Park Lots'', ''Barrie (Market House Simco Dist.)'', ''Barrie (Town lots for Jail & Court House)'', ''Barrie (Wesleyan Methodists Soc.)'', ''Barrie (Wesleyan Methodists Society)'', ''Barrie Members of the church of Scotland'', ''Barrington'', ''Barron'', ''Barry'', ''Bart'', ''Bartels'', ''Bartholomew'', ''Bartlett'', ''Bartley/Green'', ''Bartley/Lyst'', ''Bartley/Nicholson'', ''Bartley/Shewman'', ''Bartlough'', ''Bartly/Green'', ''Barton'', ''Barton/Shatford'', ''Barton/Smades'', ''Barton/Smedes'', ''Barton/Thomson'', ''Barton/Wiley'', ''Bartram'', ''Baseer'', ''Basset'', ''Bassey/Smith'', ''Bastedo'', ''Bastedo/Rogers'', ''Bastian/Stewart'', ''Bates/Ghent'', ''Bath Circuit, Wesleyan Methodist Church'', ''Bathurst'', ''Bathurst District Clergy Reserves'', ''Batter/Merkley'', ''Batter/Wright'', ''Battiam'', ''Battian'', ''Baxter'', ''Bay Indians'', ''Bay de Chaleur (Timber Agency)'', ''Bay of Quinte'', ''Bayeaux/Vincent'', ''Beach'', ''Beach Lots'', ''Beaches in Quebec'', ''Beal'', ''Beaman'', ''Beard'', ''Bearman'', ''Beasley/Stegman'', ''Beasley/Van Every'', ''Beatie'', ''Beatson'', ''Beatys'', ''BeauSoleil'', ''Beaugrand'', ''Beaumont'', ''Beaver'', ''Bebee/Smith'', ''Becker'', ''Becket'', ''Beckwith'', ''Bedal'', ''Bedell'', ''Bedford & Sheffield (Clergy Reserves & School lands)'', ''Bedford, Reservation for Schools'', ''Bedstead'', ''Beeman'', ''Beick'', ''Beikie'', ''Belding'', ''Bell'', ''Bell/Huff'', ''Bell/Scott'', ''Bell/Sharp'', ''Bell/Sills'', ''Bell/Switzer'', ''Bell/Taylor'', ''Bellamy'', ''Bellar/Le Motte'', ''Bellau'', ''Bellaw'', ''Belleperche'', ''Belleville (Moira River)'', ''Belleville town'', ''Belleville, Minister & Church Wardens of C. of E.'', ''Belleville, President and Board of Police.'', ''Belleville, Rev. Campbell for parsonage and glebe (C. of E.)'', ''Belleville, church of England, parsonage and glebe (Rev. Thomas Campbell)'', ''Bellows'', ''Beman'', ''Bencroft'', ''Bender'', ''Bender/Murry'', ''Bender/Stelmayon'', ''Bender/Wood'', ''Benedict/Lossee'', ''Benet'', ''Benett'', ''Benidick/Olds'', ''Benn/Sills'', ''Benn/Simmons'', ''Benn/Soles'', ''Bennet'', ''Bennett'', ''Bennitt'', ''Benson'', ''Benson/Richardson'', ''Benson/Steel'', ''Benson/Terell'', ''Bent'', ''Bentinck & Glenelg School Section'', ''Bentinck School Trustees'', ''Bentley/Jayne'', ''Benton'', ''Berczy'', ''Berdan/Thomson'', ''Bergeron'', ''Bernard'', ''Berney'', ''Bernier'', ''Berry'', ''Berry/Ulman'', ''Berry/Willman'', ''Bertram'', ''Besserer'', ''Bessey'', ''Bessey/Price'', ''Bessil/Michel'', ''Bessy/Newkirk/Newcark'', ''Best'', ''Bethune'', ''Bethune/McKenzie'', ''Bethune/Wilkinson'', ''Bettys'', ''Bettys/Mahon'', ''Bieges'', ''Biggars'', ''Bighard'', ''Bignal'', ''Bigsby'', ''Bill'', ''Billows'', ''Bills'', ''Birch'', ''Birch/Rowel'', ''Bird'', ''Birdsale/McDonald'', ''Birdsall'', ''Birdsall/Shotwell'', ''Birmingham'', ''Birn'', ''Bisco'', ''Bishop'', ''Bishop of Nova Scotia'', ''Bishop of Quebec'', ''Bishop of Toronto'', ''Bishop of York'', ''Bishop, Coadjutor and Administrator of Regiopolis'', ''Bissell/Landon'', ''Black'', ''Black/Sparks'', ''Black/Umphrey'', ''Blackburn'', ''Blackley/Murrison'', ''Blackman'', ''Blackwell'', ''Blackwood'', ''Blades'', ''Blain'', ''Blair'', ''Blake'', ''Blakeby/McDonald'', ''Blakely/Mix'', ''Blakely/Toppin'', ''Blaker/Ripsom'', ''Blakley/Ruscow'', ''Blan'', ''Blanchard'', ''Blandford, Brock District (Purchase of Clergy Reserve)'', ''Blasdell'', ''Bleecker'', ''Bleex'', ''Blond'', ''Blondheim'', ''Board for Education'', ''Board of Ordnance'', ''Board of Police'', ''Board of Police of town of Belleville.'', ''Board of School Trustees'', ''Board of Trustees of the Thomas County Grammar School'', ''Board of Trustees, Grammar School'', ''Board of Works'', ''Boardman'', ''Bodart'', ''Bogart'', ''Bogart/Spicer'', ''Bogert/Matheson'', ''Boice/McGuire'', ''Boice/Springer'', ''Boice/Van Wicklin'', ''Boid/Shipman'', ''Boies/Leonard'', ''Boisbouscache village plot'', ''Bolduc'', ''Bolton'', ''Bolton/Sax'', ''Bond'', ''Bone'', ''Bonesteell/Willcox'', ''Bonet'', ''Bonnycastle'', ''Bonsteel/Laushway'', ''Bonsteel/Smith'', ''Booth'', ''Booth/Stalker'', ''Booth/Terry'', ''Boreson/Manhart'', ''Borrat'', ''Borrie'', ''Borrit/Wright'', ''Boshart'', ''Bostwick'', ''Bosworth'', ''Botsford'', ''Boucher'', ''Boucher/Loucks'', ''Bouchette/Reynolds'', ''Bouchier'', ''Bouchner/Willson'', ''Bouchuer/Stringer'', ''Bouck'', ''Bouck/Garlough'', ''Bouck/Reid'', ''Boughner/Owen'', ''Boughner/Willson'', ''Boulton'', ''Boulton/McMurphy'', ''Boundary Commission'', ''Boundary Lines'', ''Boundary line the township Beverly'', ''Boundary lines of lots'', ''Bourne'', ''Boutlon/Rawley'', ''Bowac/Markle'', ''Bowdish'', ''Bowen'', ''Bowen/More'', ''Bowen/Sager'', ''Bowen/Shaver'', ''Bower/Scott'', ''Bower/Worden'', ''Bowerman'', ''Bowers/Scott'', ''Bowes'', ''Bowie'', ''Bowling/Goslee'', ''Bowman'', ''Bowman/Gisso'', ''Bowman/Schram'', ''Bowman/Scott'', ''Bowman/Snure'''

# Row 7 (title) updates
$ws.Range("P7").Value = 2
$ws.Range("R7").Value = 2247940
$ws.Range("S7").Value = 2
$ws.Range("T7").Value = '[''0'', ''Upper Canada Land Petitions'']'

# Row 8 (reference) updates
$ws.Range("P8").Value = 2
$ws.Range("R8").Value = 807275
$ws.Range("S8").Value = 3
$ws.Range("T8").Value = '[''0'', ''RG 1 L3'', ''RG 5 A1, Land Petitions in Upper Canada Sundries'']'
